# Applies the OOXML edits to the VML screenshot shapes and the title-block
# text boxes. This document's drawings/text boxes are legacy w:pict/VML
# content, which Word's Shapes/Find APIs don't reach into here, so the
# fix is applied by feeding the corrected body XML back through
# Range.InsertXML (Document.Content.InsertXML replaces the whole body).
$d = $word.ActiveDocument
$bodyXml = '<w:p w14:paraId="79FEA289" w14:textId="77777777" w:rsidR="00B827FD" w:rsidRDefault="00B827FD" w:rsidP="002F2B70"><w:pPr><w:pStyle w:val="a1"/><w:jc w:val="left"/></w:pPr></w:p><w:p w14:paraId="06730B81" w14:textId="77777777" w:rsidR="00373697" w:rsidRDefault="00373697" w:rsidP="00B827FD"><w:pPr><w:pStyle w:val="a1"/><w:jc w:val="center"/></w:pPr></w:p><w:p w14:paraId="00BCFC07" w14:textId="77777777" w:rsidR="00373697" w:rsidRDefault="00373697" w:rsidP="00373697"><w:pPr><w:pStyle w:val="a1"/><w:jc w:val="center"/></w:pPr><w:r><w:rPr><w:lang w:val="ru-RU"/></w:rPr><w:t xml:space="preserve">     </w:t></w:r><w:r><w:pict w14:anchorId="6B9323F2"><v:shapetype id="_x0000_t75" coordsize="21600,21600" o:spt="75" o:preferrelative="t" path="m@4@5l@4@11@9@11@9@5xe" filled="f" stroked="f"><v:stroke joinstyle="miter"/><v:formulas><v:f eqn="if lineDrawn pixelLineWidth 0"/><v:f eqn="sum @0 1 0"/><v:f eqn="sum 0 0 @1"/><v:f eqn="prod @2 1 2"/><v:f eqn="prod @3 21600 pixelWidth"/><v:f eqn="prod @3 21600 pixelHeight"/><v:f eqn="sum @0 0 1"/><v:f eqn="prod @6 1 2"/><v:f eqn="prod @7 21600 pixelWidth"/><v:f eqn="sum @8 21600 0"/><v:f eqn="prod @7 21600 pixelHeight"/><v:f eqn="sum @10 21600 0"/></v:formulas><v:path o:extrusionok="f" gradientshapeok="t" o:connecttype="rect"/><o:lock v:ext="edit" aspectratio="t"/></v:shapetype><v:shape id="_x0000_i1025" type="#_x0000_t75" style="width:495.4pt;height:241.3pt" o:bordertopcolor="this" o:borderleftcolor="this" o:borderbottomcolor="this" o:borderrightcolor="this"><v:imagedata r:id="rId7" o:title=""/><w10:bordertop type="single" width="4"/><w10:borderleft type="single" width="4"/><w10:borderbottom type="single" width="4"/><w10:borderright type="single" width="4"/></v:shape></w:pict></w:r></w:p><w:p w14:paraId="6D1D523C" w14:textId="77777777" w:rsidR="00373697" w:rsidRDefault="00373697" w:rsidP="00373697"><w:pPr><w:pStyle w:val="a1"/><w:jc w:val="center"/></w:pPr></w:p><w:p w14:paraId="20120A29" w14:textId="18E8C2C2" w:rsidR="00627F81" w:rsidRPr="00B827FD" w:rsidRDefault="00373697" w:rsidP="00373697"><w:pPr><w:pStyle w:val="a1"/><w:jc w:val="center"/></w:pPr><w:r><w:rPr><w:lang w:val="ru-RU"/></w:rPr><w:t xml:space="preserve">     </w:t></w:r><w:r><w:pict w14:anchorId="06D6FD38"><v:shape id="_x0000_i1026" type="#_x0000_t75" style="width:496.75pt;height:236.85pt" o:bordertopcolor="this" o:borderleftcolor="this" o:borderbottomcolor="this" o:borderrightcolor="this"><v:imagedata r:id="rId8" o:title=""/><w10:bordertop type="single" width="4"/><w10:borderleft type="single" width="4"/><w10:borderbottom type="single" width="4"/><w10:borderright type="single" width="4"/></v:shape></w:pict></w:r><w:r><w:br w:type="page"/></w:r><w:r w:rsidR="00B827FD"><w:rPr><w:noProof/></w:rPr><w:lastRenderedPageBreak/><w:pict w14:anchorId="5CA24E01"><v:group id="Группа 1293" o:spid="_x0000_s1961" style="position:absolute;left:0;text-align:left;margin-left:11pt;margin-top:-4.1pt;width:522.15pt;height:791.7pt;z-index:1;mso-position-horizontal-relative:margin;mso-position-vertical-relative:margin" coordorigin="1124,774" coordsize="10441,15504" o:gfxdata="UEsDBBQABgAIAAAAIQC2gziS/gAAAOEBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbJSRQU7DMBBF&#xA;90jcwfIWJU67QAgl6YK0S0CoHGBkTxKLZGx5TGhvj5O2G0SRWNoz/78nu9wcxkFMGNg6quQqL6RA&#xA;0s5Y6ir5vt9lD1JwBDIwOMJKHpHlpr69KfdHjyxSmriSfYz+USnWPY7AufNIadK6MEJMx9ApD/oD&#xA;OlTrorhX2lFEilmcO2RdNtjC5xDF9pCuTyYBB5bi6bQ4syoJ3g9WQ0ymaiLzg5KdCXlKLjvcW893&#xA;SUOqXwnz5DrgnHtJTxOsQfEKIT7DmDSUCaxw7Rqn8787ZsmRM9e2VmPeBN4uqYvTtW7jvijg9N/y&#xA;JsXecLq0q+WD6m8AAAD//wMAUEsDBBQABgAIAAAAIQA4/SH/1gAAAJQBAAALAAAAX3JlbHMvLnJl&#xA;bHOkkMFqwzAMhu+DvYPRfXGawxijTi+j0GvpHsDYimMaW0Yy2fr2M4PBMnrbUb/Q94l/f/hMi1qR&#xA;JVI2sOt6UJgd+ZiDgffL8ekFlFSbvV0oo4EbChzGx4f9GRdb25HMsYhqlCwG5lrLq9biZkxWOiqY&#xA;22YiTra2kYMu1l1tQD30/bPm3wwYN0x18gb45AdQl1tp5j/sFB2T0FQ7R0nTNEV3j6o9feQzro1i&#xA;OWA14Fm+Q8a1a8+Bvu/d/dMb2JY5uiPbhG/ktn4cqGU/er3pcvwCAAD//wMAUEsDBBQABgAIAAAA&#xA;IQDr+HkYhgsAAJWkAAAOAAAAZHJzL2Uyb0RvYy54bWzsXWtu48gR/h8gdyD43yO+RQqjWcza1iDA&#xA;JFlkNwegJUoiQpEKSY88uwgQYI+wF8kNcoXdG6W6utlqkvLass3O0FMeYCyaEtVsVn/11bPffnO3&#xA;y4xPSVmlRT437TeWaST5slil+WZu/v2HxUVoGlUd56s4K/Jkbn5OKvObd3/8w9vDfpY4xbbIVklp&#xA;wEXyanbYz81tXe9nk0m13Ca7uHpT7JMcTq6LchfXcFhuJqsyPsDVd9nEsaxgcijK1b4slklVwV+v&#xA;+EnzHV5/vU6W9V/X6yqpjWxuwthq/L/E/2/Y/5N3b+PZpoz323QphhE/YRS7OM3hS+WlruI6Nm7L&#xA;tHepXbosi6pY12+WxW5SrNfpMsF7gLuxrc7dfCiL2z3ey2Z22OzlNMHUdubpyZdd/uXTd6WRruDZ&#xA;OZFrGnm8g6f06y+//fu3n3/9L/z7j4EnYJ4O+80M3v6h3H+//67kNwsvPxbLf1RwetI9z443/M3G&#xA;zeHPxQouHN/WBc7T3brcsUvADBh3+Dg+y8eR3NXGEv4YBK7tWr5pLOGcbVm+50fiiS238FjZB23b&#xA;8UwDzk+nHn+Yy+21+LxteZ4tPu37Fp6fxDP+1ThcMTx2byB+1XGGq+fN8PfbeJ/gg6vYlB1nGMbK&#xA;Z/hvIJlxvskSI3JsNnA2AnhrM7UVn1cjLy638L7kfVkWh20Sr2Bk+H4Yv/IBdlDBU3lwovvz1cz2&#xA;785WPNuXVf0hKXYGezE3Sxg/Psj408eqZo//+Bb2XPNikWYZ/D2eZblxmJuO71kWfqIqsnTFzrKT&#xA;Vbm5ucxK41PMFij+sNmAq7XetktrgIks3c3NUL4pnrEJuc5X+DV1nGb8NXw4y9nF4dZgcOIVX44/&#xA;RVZ0HV6H3oXnBNcXnnV1dfF+celdBAt76l+5V5eXV/a/2Dhtb7ZNV6skZ0NtoMH2HicYAqT4opbg&#xA;0Lql1p0v8Kd/55P2MHBi4K6a33h3KAns4XMhuilWn0EQyoJjHWAzvNgW5Y+mcQCcm5vVP2/jMjGN&#xA;7E85CFNkex4DRjzw/KkDB6V65kY9E+dLuNTcrE2Dv7ysOZje7st0s4VvsvEZ58V7WOnrFCWDCScf&#xA;FYybHcBK07fkAD/4kvuY5my1Ocpqu8w5kC3vcgFkcsHhEv7h8x5Aq7Xe+Ecev94Cz0F8sj0rEggl&#xA;VxzHphAgjot8A4qdtZbBwHFW9aw1UDpiSY1seUnMefKqYY9Bu4AGHQF19QooqPdGQMMp+2rETVTA&#xA;oBKYdmXqF/WuVJ5HrBfqgES0y+lOa4CRiui0I6IIZGyhAGMZHkMdJwAQZxzQsy0EY0VEuXg6dhAS&#xA;iKIVBYzhORxlpBIKZmZLy6NK1SahbhA1Wp4kdGgWPVIJBS3bktBAq5r3fA9AnDBUi503Sgl1wTxv&#xA;SyiyQW0Y6ts2gDhK6H2WEmn5l/JEjFRCu6YScj5tEhq5UYOhkY1W2pGHguuUmUm+i4uGDKXZ10hD&#xA;XbsLoZFWJQ/OXWHL+1NXeMobZxPILpxDT/oLG0q2M32GZ5e8TTrdoa4NARGVhnIx0QahRwkNLB6Q&#xA;USBUSuhLO5tIQlvBjCZ60Py+P4qg3x/qYACCS6gSI3P1xsj8xpoPfETwo5B6gYgoOhA1+V2v/VkR&#xA;MuZ1VTjZ649f8RnFCPlPEcRwrW+d6GIRhNMLb+H5F9HUCi8sO/o2Ciwv8q4W7ZAchnJ40sAzeAYL&#xA;S7oQ9cPYyv2xuReLSspFyIbfLL7m96lFWN/d3IGMMXQ+M6rHEU9E9cQBj+qJAx7VEwdji+oBSEg1&#xA;poKEGtobPJAeTIHzo7naAwkfchGQahFIqOH9s4PcBBIsEeKBeD8DCZ7A09gahBVK0g1gBSjzPqFQ&#xA;o6xDY4XjWo3joIcVtusRWJzIBSKweESe05MYBYKFVJUEFm2wkDmQKrFQ491Dg4U7dZpodw8sQgh/&#xA;E7FgxtLzsueIWJxDLKSqJKxoYwUo7j6xUDMPhsYKlq95rxHCeA/z95IRQkbImfnTT+cVUlMSVrSx&#xA;QqYhq7xCzQEZGisizPNGh0XkYOCUa0FWaTGdQuAKscJFf+f98UvyakLZwr35xcQrzuEVUlMSVrSx&#xA;QqY5qFihZuPowwrfgYhhKyn8iBUUAVFql8hfMai/QmpKwoo2VsjUfBUr1LyoobHCZxE0HghxPIT0&#xA;I68IXA9y+jAxykaPExELLHQksBgULKSqJLBog0WnSoKTfY3JP5BtxJEisjoWiFJqxsstmzptKjVT&#xA;+yGcrPN/PaVmEKzrlElwjqlTRGWkzu0os+EyKJ9XG08ZlBozKB2WWNvKoIQWEyJHR0sx5DGD0nMD&#xA;yvF1B23YoKTkSZco71zxcF4G5x4sRZChl76WByChMjkK8+Qi6EDw/5FQaNiiDUN5KhuGBp/QX4Qw&#xA;VK+EypQcIaEyxKgZQ307vBdDH8jvzc7tykES+jIYKrpKsbxYbEgluisB7sncDWycZUDEgwPfcH2r&#xA;bOie05g0Xb7o+CzrS4nAxjPZtorbQUL0lI5VEfyIzzALmQ1fukxAj8iGWpra58CcnopxczePIOVD&#xA;+5daU9VUM4WhD05ynNvePPUsRopZUcwqnqn9wSQQnZOJj3lz0rNKrqWWawnAr58L40mvvYbOdhDU&#xA;Bt4L8e2OQretKSssJ6zo99sjN/SgbmhZ/DsWrDgyDE6usD3mCZ4l49MNzxIOdx08y3e6daMKz+IO&#xA;7VHyrFNxPAieHQ134lnj7thJuUHn5AZBJwUh+mMBT+Hlk02aB2khDAaZjOEpAX9PqhoiWq+ivS+B&#xA;xVlgIUMwYwGLxzItGQ0VTMsXsKiFaXnQMQZA+JjO8yqYFut23ava8KUEaQBQ8mixDRyG641O4HkW&#xA;eMro4FjAUxPTYntH9IFCzpYGoCCXVrPbC4HFS26P8XT3N297yBbgWMDikUyL7dLAF3vDtERIWgvT&#xA;mnJP+UmmxfOLxujT8mQ8VjFVfVkhqAFAiWkR03rpvYWeAZ5S9McCnrqY1qkkA546xkagASiIaRHT&#xA;Uor5MCzYZN+0OgNqK6S3R1cd+1imJTMFGqYlsgR0MK3AshGETzMtkTkwviwtoI8nTFVZBqUBQIlp&#xA;EdP6gpiWTDwiptVK02KbwPR9WjrTDIhpEdP6wpiW1JNjAYs209Kwc7EzZXX8HDh4xQbviy6ss+G3&#xA;AAx95o3nzYu6W6/ACZbcCY1UH+hddHbFBtVlSmfHl1/1Bk+/r9p4haQmJwIUE8EYmCj2u+a4Tgh2&#xA;D8qpBWUhal0H1SssFmJCFFObemwp+5nLZXh2vQKkvlMeHds9W2CALBaDTJqWQoN9PcREgbWsRaGB&#xA;RkWssEMMNR1dEq4FzXIQKxoXFG0M/tXtaetMISepLaFq6oEWCRWtMKCMG+3pUxKKq0aWKPZUGVGu&#xA;17svOEioDJpzoyDQXMZtOX6zB0tvQ0aBoLQhYyYy0p+xUVJndyrJRr58o4DJR8/fFcg4rAbHeOiH&#xA;AkY9cH0xjnGE0WnQFHtT213FJ0RliYOWJdqyWGQs/i7G3aGBztClNZCadAIsZBxWA1jYFu7Czf1c&#xA;XbSwHYgJ8zpmgguCi7ROSiNLd3MzHBQuHFFyMp6UT11wIRMZlPTGQEYeNcFFY6M5VseLQHCx4vyU&#xA;NgsCyqktbwl25iSH40mHo8zREcayDDvqcTiGTgMVdtixQ0QEzQnQw/SC3hxqeTYmWxnkQw3xBjJ0&#xA;oEVAI+hsKj3incJTEtDVKsmZOvuqnTkyB0ElXNKa1UC4Qj8C0odRm4CXDB+dOQ5z76B55kEjKtCB&#xA;9+NomSxr3HQ8/gTbfvO37suq/pAUu5avLRYbwSt9aGlreNoaHtwhctdn2OWP+NZJvnUqHwT2OROz&#xA;pQMtglA0B+03rfJZGwKWDuJCP3ECi0zkDpPrd1hfjgzQjcX1i6mOh80e1emmjPfbdHkV17F6DK8P&#xA;+1niFNsiWyXlu/8BAAD//wMAUEsDBBQABgAIAAAAIQABKZEw4gAAAAwBAAAPAAAAZHJzL2Rvd25y&#xA;ZXYueG1sTI/BSsNAEIbvgu+wjOCt3aSarY3ZlFLUUxFsBfE2TaZJaHY3ZLdJ+vaOJ739w3z88022&#xA;nkwrBup946yGeB6BIFu4srGVhs/D6+wJhA9oS2ydJQ1X8rDOb28yTEs32g8a9qESXGJ9ihrqELpU&#xA;Sl/UZNDPXUeWdyfXGww89pUsexy53LRyEUVKGmwsX6ixo21NxXl/MRreRhw3D/HLsDufttfvQ/L+&#xA;tYtJ6/u7afMMItAU/mD41Wd1yNnp6C629KLVMFuqR0Y5LNQKBBMrteRwZDRJVAIyz+T/J/IfAAAA&#xA;//8DAFBLAQItABQABgAIAAAAIQC2gziS/gAAAOEBAAATAAAAAAAAAAAAAAAAAAAAAABbQ29udGVu&#xA;dF9UeXBlc10ueG1sUEsBAi0AFAAGAAgAAAAhADj9If/WAAAAlAEAAAsAAAAAAAAAAAAAAAAALwEA&#xA;AF9yZWxzLy5yZWxzUEsBAi0AFAAGAAgAAAAhAOv4eRiGCwAAlaQAAA4AAAAAAAAAAAAAAAAALgIA&#xA;AGRycy9lMm9Eb2MueG1sUEsBAi0AFAAGAAgAAAAhAAEpkTDiAAAADAEAAA8AAAAAAAAAAAAAAAAA&#xA;4A0AAGRycy9kb3ducmV2LnhtbFBLBQYAAAAABAAEAPMAAADvDgAAAAA=&#xA;"><v:rect id="Rectangle 921" o:spid="_x0000_s1962" style="position:absolute;left:1124;top:774;width:10441;height:15504;visibility:visible" o:gfxdata="UEsDBBQABgAIAAAAIQDb4fbL7gAAAIUBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbHyQz07DMAyH&#xA;70i8Q+QralM4IITa7kDhCAiNB7ASt43WOlEcyvb2pNu4IODoPz9/n1xv9vOkForiPDdwXVagiI23&#xA;jocG3rdPxR0oScgWJ8/UwIEENu3lRb09BBKV0ywNjCmFe63FjDSjlD4Q50nv44wpl3HQAc0OB9I3&#xA;VXWrjedEnIq03oC27qjHjympx31un0wiTQLq4bS4shrAECZnMGVTvbD9QSnOhDInjzsyuiBXWQP0&#xA;r4R18jfgnHvJr4nOknrFmJ5xzhraRtHWf3Kkpfz/yGo5S+H73hkquyhdjr3R8m2lj09svwAAAP//&#xA;AwBQSwMEFAAGAAgAAAAhAFr0LFu/AAAAFQEAAAsAAABfcmVscy8ucmVsc2zPwWrDMAwG4Ptg72B0&#xA;X5TuUMaI01uh19I+gLGVxCy2jGSy9e1nemrHjpL4P0nD4SetZiPRyNnCruvBUPYcYp4tXC/Htw8w&#xA;Wl0ObuVMFm6kcBhfX4Yzra62kC6xqGlKVgtLreUTUf1CyWnHhXKbTCzJ1VbKjMX5LzcTvvf9HuXR&#xA;gPHJNKdgQU5hB+ZyK23zHztFL6w81c5zQp6m6P9TMfB3PtPWFCczVQtB9N4U2rp2HOA44NMz4y8A&#xA;AAD//wMAUEsDBBQABgAIAAAAIQDU3msFwgAAAN0AAAAPAAAAZHJzL2Rvd25yZXYueG1sRE/NisIw&#xA;EL4L+w5hFrxpuiJiq1GqIHiS3eoDDM1sW2wmtYlt9ek3woK3+fh+Z70dTC06al1lWcHXNAJBnFtd&#xA;caHgcj5MliCcR9ZYWyYFD3Kw3XyM1pho2/MPdZkvRAhhl6CC0vsmkdLlJRl0U9sQB+7XtgZ9gG0h&#xA;dYt9CDe1nEXRQhqsODSU2NC+pPya3Y2Cqx+6U1pkz0N82cX59y7t77dUqfHnkK5AeBr8W/zvPuow&#xA;fxbP4fVNOEFu/gAAAP//AwBQSwECLQAUAAYACAAAACEA2+H2y+4AAACFAQAAEwAAAAAAAAAAAAAA&#xA;AAAAAAAAW0NvbnRlbnRfVHlwZXNdLnhtbFBLAQItABQABgAIAAAAIQBa9CxbvwAAABUBAAALAAAA&#xA;AAAAAAAAAAAAAB8BAABfcmVscy8ucmVsc1BLAQItABQABgAIAAAAIQDU3msFwgAAAN0AAAAPAAAA&#xA;AAAAAAAAAAAAAAcCAABkcnMvZG93bnJldi54bWxQSwUGAAAAAAMAAwC3AAAA9gIAAAAA&#xA;" filled="f" strokeweight="2pt"/><v:line id="Line 922" o:spid="_x0000_s1963" style="position:absolute;visibility:visible" from="1642,14094" to="1643,14899" o:connectortype="straight" o:gfxdata="UEsDBBQABgAIAAAAIQDb4fbL7gAAAIUBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbHyQz07DMAyH&#xA;70i8Q+QralM4IITa7kDhCAiNB7ASt43WOlEcyvb2pNu4IODoPz9/n1xv9vOkForiPDdwXVagiI23&#xA;jocG3rdPxR0oScgWJ8/UwIEENu3lRb09BBKV0ywNjCmFe63FjDSjlD4Q50nv44wpl3HQAc0OB9I3&#xA;VXWrjedEnIq03oC27qjHjympx31un0wiTQLq4bS4shrAECZnMGVTvbD9QSnOhDInjzsyuiBXWQP0&#xA;r4R18jfgnHvJr4nOknrFmJ5xzhraRtHWf3Kkpfz/yGo5S+H73hkquyhdjr3R8m2lj09svwAAAP//&#xA;AwBQSwMEFAAGAAgAAAAhAFr0LFu/AAAAFQEAAAsAAABfcmVscy8ucmVsc2zPwWrDMAwG4Ptg72B0&#xA;X5TuUMaI01uh19I+gLGVxCy2jGSy9e1nemrHjpL4P0nD4SetZiPRyNnCruvBUPYcYp4tXC/Htw8w&#xA;Wl0ObuVMFm6kcBhfX4Yzra62kC6xqGlKVgtLreUTUf1CyWnHhXKbTCzJ1VbKjMX5LzcTvvf9HuXR&#xA;gPHJNKdgQU5hB+ZyK23zHztFL6w81c5zQp6m6P9TMfB3PtPWFCczVQtB9N4U2rp2HOA44NMz4y8A&#xA;AAD//wMAUEsDBBQABgAIAAAAIQDUoGg2vgAAAN0AAAAPAAAAZHJzL2Rvd25yZXYueG1sRE+9CsIw&#xA;EN4F3yGc4KapgqLVKCJU3MTq4nY2Z1tsLqWJWt/eCILbfXy/t1y3phJPalxpWcFoGIEgzqwuOVdw&#xA;PiWDGQjnkTVWlknBmxysV93OEmNtX3ykZ+pzEULYxaig8L6OpXRZQQbd0NbEgbvZxqAPsMmlbvAV&#xA;wk0lx1E0lQZLDg0F1rQtKLunD6PgfjlPkt1hq09VutHXPPGX600r1e+1mwUIT63/i3/uvQ7zx/MJ&#xA;fL8JJ8jVBwAA//8DAFBLAQItABQABgAIAAAAIQDb4fbL7gAAAIUBAAATAAAAAAAAAAAAAAAAAAAA&#xA;AABbQ29udGVudF9UeXBlc10ueG1sUEsBAi0AFAAGAAgAAAAhAFr0LFu/AAAAFQEAAAsAAAAAAAAA&#xA;AAAAAAAAHwEAAF9yZWxzLy5yZWxzUEsBAi0AFAAGAAgAAAAhANSgaDa+AAAA3QAAAA8AAAAAAAAA&#xA;AAAAAAAABwIAAGRycy9kb3ducmV2LnhtbFBLBQYAAAAAAwADALcAAADyAgAAAAA=&#xA;" strokeweight="2pt"/><v:line id="Line 923" o:spid="_x0000_s1964" style="position:absolute;visibility:visible" from="1129,14087" to="11553,14087" o:connectortype="straight" o:gfxdata="UEsDBBQABgAIAAAAIQDb4fbL7gAAAIUBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbHyQz07DMAyH&#xA;70i8Q+QralM4IITa7kDhCAiNB7ASt43WOlEcyvb2pNu4IODoPz9/n1xv9vOkForiPDdwXVagiI23&#xA;jocG3rdPxR0oScgWJ8/UwIEENu3lRb09BBKV0ywNjCmFe63FjDSjlD4Q50nv44wpl3HQAc0OB9I3&#xA;VXWrjedEnIq03oC27qjHjympx31un0wiTQLq4bS4shrAECZnMGVTvbD9QSnOhDInjzsyuiBXWQP0&#xA;r4R18jfgnHvJr4nOknrFmJ5xzhraRtHWf3Kkpfz/yGo5S+H73hkquyhdjr3R8m2lj09svwAAAP//&#xA;AwBQSwMEFAAGAAgAAAAhAFr0LFu/AAAAFQEAAAsAAABfcmVscy8ucmVsc2zPwWrDMAwG4Ptg72B0&#xA;X5TuUMaI01uh19I+gLGVxCy2jGSy9e1nemrHjpL4P0nD4SetZiPRyNnCruvBUPYcYp4tXC/Htw8w&#xA;Wl0ObuVMFm6kcBhfX4Yzra62kC6xqGlKVgtLreUTUf1CyWnHhXKbTCzJ1VbKjMX5LzcTvvf9HuXR&#xA;gPHJNKdgQU5hB+ZyK23zHztFL6w81c5zQp6m6P9TMfB3PtPWFCczVQtB9N4U2rp2HOA44NMz4y8A&#xA;AAD//wMAUEsDBBQABgAIAAAAIQAkcvZBvgAAAN0AAAAPAAAAZHJzL2Rvd25yZXYueG1sRE+9CsIw&#xA;EN4F3yGc4KapgqLVKCJU3MTq4nY2Z1tsLqWJWt/eCILbfXy/t1y3phJPalxpWcFoGIEgzqwuOVdw&#xA;PiWDGQjnkTVWlknBmxysV93OEmNtX3ykZ+pzEULYxaig8L6OpXRZQQbd0NbEgbvZxqAPsMmlbvAV&#xA;wk0lx1E0lQZLDg0F1rQtKLunD6PgfjlPkt1hq09VutHXPPGX600r1e+1mwUIT63/i3/uvQ7zx/Mp&#xA;fL8JJ8jVBwAA//8DAFBLAQItABQABgAIAAAAIQDb4fbL7gAAAIUBAAATAAAAAAAAAAAAAAAAAAAA&#xA;AABbQ29udGVudF9UeXBlc10ueG1sUEsBAi0AFAAGAAgAAAAhAFr0LFu/AAAAFQEAAAsAAAAAAAAA&#xA;AAAAAAAAHwEAAF9yZWxzLy5yZWxzUEsBAi0AFAAGAAgAAAAhACRy9kG+AAAA3QAAAA8AAAAAAAAA&#xA;AAAAAAAABwIAAGRycy9kb3ducmV2LnhtbFBLBQYAAAAAAwADALcAAADyAgAAAAA=&#xA;" strokeweight="2pt"/><v:line id="Line 924" o:spid="_x0000_s1965" style="position:absolute;visibility:visible" from="2265,14101" to="2266,16269" o:connectortype="straight" o:gfxdata="UEsDBBQABgAIAAAAIQDb4fbL7gAAAIUBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbHyQz07DMAyH&#xA;70i8Q+QralM4IITa7kDhCAiNB7ASt43WOlEcyvb2pNu4IODoPz9/n1xv9vOkForiPDdwXVagiI23&#xA;jocG3rdPxR0oScgWJ8/UwIEENu3lRb09BBKV0ywNjCmFe63FjDSjlD4Q50nv44wpl3HQAc0OB9I3&#xA;VXWrjedEnIq03oC27qjHjympx31un0wiTQLq4bS4shrAECZnMGVTvbD9QSnOhDInjzsyuiBXWQP0&#xA;r4R18jfgnHvJr4nOknrFmJ5xzhraRtHWf3Kkpfz/yGo5S+H73hkquyhdjr3R8m2lj09svwAAAP//&#xA;AwBQSwMEFAAGAAgAAAAhAFr0LFu/AAAAFQEAAAsAAABfcmVscy8ucmVsc2zPwWrDMAwG4Ptg72B0&#xA;X5TuUMaI01uh19I+gLGVxCy2jGSy9e1nemrHjpL4P0nD4SetZiPRyNnCruvBUPYcYp4tXC/Htw8w&#xA;Wl0ObuVMFm6kcBhfX4Yzra62kC6xqGlKVgtLreUTUf1CyWnHhXKbTCzJ1VbKjMX5LzcTvvf9HuXR&#xA;gPHJNKdgQU5hB+ZyK23zHztFL6w81c5zQp6m6P9TMfB3PtPWFCczVQtB9N4U2rp2HOA44NMz4y8A&#xA;AAD//wMAUEsDBBQABgAIAAAAIQBLPlPawgAAAN0AAAAPAAAAZHJzL2Rvd25yZXYueG1sRE9Li8Iw&#xA;EL4L+x/CLHjTdAVX7TaKCBVvsq0Xb2MzfWAzKU3U+u/NwoK3+fiek2wG04o79a6xrOBrGoEgLqxu&#xA;uFJwytPJEoTzyBpby6TgSQ42649RgrG2D/6le+YrEULYxaig9r6LpXRFTQbd1HbEgSttb9AH2FdS&#xA;9/gI4aaVsyj6lgYbDg01drSrqbhmN6Pgej7N0/1xp/M22+pLlfrzpdRKjT+H7Q8IT4N/i//dBx3m&#xA;z1YL+PsmnCDXLwAAAP//AwBQSwECLQAUAAYACAAAACEA2+H2y+4AAACFAQAAEwAAAAAAAAAAAAAA&#xA;AAAAAAAAW0NvbnRlbnRfVHlwZXNdLnhtbFBLAQItABQABgAIAAAAIQBa9CxbvwAAABUBAAALAAAA&#xA;AAAAAAAAAAAAAB8BAABfcmVscy8ucmVsc1BLAQItABQABgAIAAAAIQBLPlPawgAAAN0AAAAPAAAA&#xA;AAAAAAAAAAAAAAcCAABkcnMvZG93bnJldi54bWxQSwUGAAAAAAMAAwC3AAAA9gIAAAAA&#xA;" strokeweight="2pt"/><v:line id="Line 925" o:spid="_x0000_s1966" style="position:absolute;visibility:visible" from="3692,14101" to="3693,16269" o:connectortype="straight" o:gfxdata="UEsDBBQABgAIAAAAIQDb4fbL7gAAAIUBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbHyQz07DMAyH&#xA;70i8Q+QralM4IITa7kDhCAiNB7ASt43WOlEcyvb2pNu4IODoPz9/n1xv9vOkForiPDdwXVagiI23&#xA;jocG3rdPxR0oScgWJ8/UwIEENu3lRb09BBKV0ywNjCmFe63FjDSjlD4Q50nv44wpl3HQAc0OB9I3&#xA;VXWrjedEnIq03oC27qjHjympx31un0wiTQLq4bS4shrAECZnMGVTvbD9QSnOhDInjzsyuiBXWQP0&#xA;r4R18jfgnHvJr4nOknrFmJ5xzhraRtHWf3Kkpfz/yGo5S+H73hkquyhdjr3R8m2lj09svwAAAP//&#xA;AwBQSwMEFAAGAAgAAAAhAFr0LFu/AAAAFQEAAAsAAABfcmVscy8ucmVsc2zPwWrDMAwG4Ptg72B0&#xA;X5TuUMaI01uh19I+gLGVxCy2jGSy9e1nemrHjpL4P0nD4SetZiPRyNnCruvBUPYcYp4tXC/Htw8w&#xA;Wl0ObuVMFm6kcBhfX4Yzra62kC6xqGlKVgtLreUTUf1CyWnHhXKbTCzJ1VbKjMX5LzcTvvf9HuXR&#xA;gPHJNKdgQU5hB+ZyK23zHztFL6w81c5zQp6m6P9TMfB3PtPWFCczVQtB9N4U2rp2HOA44NMz4y8A&#xA;AAD//wMAUEsDBBQABgAIAAAAIQA6oceoxAAAAN0AAAAPAAAAZHJzL2Rvd25yZXYueG1sRI9Bi8JA&#xA;DIXvgv9hiLA3nSqs7HYdRYSKt8Xai7fYiW2xkymdUeu/N4eFvSW8l/e+rDaDa9WD+tB4NjCfJaCI&#xA;S28brgwUp2z6BSpEZIutZzLwogCb9Xi0wtT6Jx/pkcdKSQiHFA3UMXap1qGsyWGY+Y5YtKvvHUZZ&#xA;+0rbHp8S7lq9SJKldtiwNNTY0a6m8pbfnYHbufjM9r87e2rzrb1UWTxfrtaYj8mw/QEVaYj/5r/r&#xA;gxX8xbfgyjcygl6/AQAA//8DAFBLAQItABQABgAIAAAAIQDb4fbL7gAAAIUBAAATAAAAAAAAAAAA&#xA;AAAAAAAAAABbQ29udGVudF9UeXBlc10ueG1sUEsBAi0AFAAGAAgAAAAhAFr0LFu/AAAAFQEAAAsA&#xA;AAAAAAAAAAAAAAAAHwEAAF9yZWxzLy5yZWxzUEsBAi0AFAAGAAgAAAAhADqhx6jEAAAA3QAAAA8A&#xA;AAAAAAAAAAAAAAAABwIAAGRycy9kb3ducmV2LnhtbFBLBQYAAAAAAwADALcAAAD4AgAAAAA=&#xA;" strokeweight="2pt"/><v:line id="Line 926" o:spid="_x0000_s1967" style="position:absolute;visibility:visible" from="4547,14101" to="4548,16269" o:connectortype="straight" o:gfxdata="UEsDBBQABgAIAAAAIQDb4fbL7gAAAIUBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbHyQz07DMAyH&#xA;70i8Q+QralM4IITa7kDhCAiNB7ASt43WOlEcyvb2pNu4IODoPz9/n1xv9vOkForiPDdwXVagiI23&#xA;jocG3rdPxR0oScgWJ8/UwIEENu3lRb09BBKV0ywNjCmFe63FjDSjlD4Q50nv44wpl3HQAc0OB9I3&#xA;VXWrjedEnIq03oC27qjHjympx31un0wiTQLq4bS4shrAECZnMGVTvbD9QSnOhDInjzsyuiBXWQP0&#xA;r4R18jfgnHvJr4nOknrFmJ5xzhraRtHWf3Kkpfz/yGo5S+H73hkquyhdjr3R8m2lj09svwAAAP//&#xA;AwBQSwMEFAAGAAgAAAAhAFr0LFu/AAAAFQEAAAsAAABfcmVscy8ucmVsc2zPwWrDMAwG4Ptg72B0&#xA;X5TuUMaI01uh19I+gLGVxCy2jGSy9e1nemrHjpL4P0nD4SetZiPRyNnCruvBUPYcYp4tXC/Htw8w&#xA;Wl0ObuVMFm6kcBhfX4Yzra62kC6xqGlKVgtLreUTUf1CyWnHhXKbTCzJ1VbKjMX5LzcTvvf9HuXR&#xA;gPHJNKdgQU5hB+ZyK23zHztFL6w81c5zQp6m6P9TMfB3PtPWFCczVQtB9N4U2rp2HOA44NMz4y8A&#xA;AAD//wMAUEsDBBQABgAIAAAAIQBV7WIzvgAAAN0AAAAPAAAAZHJzL2Rvd25yZXYueG1sRE+9CsIw&#xA;EN4F3yGc4KapgqLVKCJU3MTq4nY2Z1tsLqWJWt/eCILbfXy/t1y3phJPalxpWcFoGIEgzqwuOVdw&#xA;PiWDGQjnkTVWlknBmxysV93OEmNtX3ykZ+pzEULYxaig8L6OpXRZQQbd0NbEgbvZxqAPsMmlbvAV&#xA;wk0lx1E0lQZLDg0F1rQtKLunD6PgfjlPkt1hq09VutHXPPGX600r1e+1mwUIT63/i3/uvQ7zx/M5&#xA;fL8JJ8jVBwAA//8DAFBLAQItABQABgAIAAAAIQDb4fbL7gAAAIUBAAATAAAAAAAAAAAAAAAAAAAA&#xA;AABbQ29udGVudF9UeXBlc10ueG1sUEsBAi0AFAAGAAgAAAAhAFr0LFu/AAAAFQEAAAsAAAAAAAAA&#xA;AAAAAAAAHwEAAF9yZWxzLy5yZWxzUEsBAi0AFAAGAAgAAAAhAFXtYjO+AAAA3QAAAA8AAAAAAAAA&#xA;AAAAAAAABwIAAGRycy9kb3ducmV2LnhtbFBLBQYAAAAAAwADALcAAADyAgAAAAA=&#xA;" strokeweight="2pt"/><v:line id="Line 927" o:spid="_x0000_s1968" style="position:absolute;visibility:visible" from="5118,14094" to="5119,16262" o:connectortype="straight" o:gfxdata="UEsDBBQABgAIAAAAIQDb4fbL7gAAAIUBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbHyQz07DMAyH&#xA;70i8Q+QralM4IITa7kDhCAiNB7ASt43WOlEcyvb2pNu4IODoPz9/n1xv9vOkForiPDdwXVagiI23&#xA;jocG3rdPxR0oScgWJ8/UwIEENu3lRb09BBKV0ywNjCmFe63FjDSjlD4Q50nv44wpl3HQAc0OB9I3&#xA;VXWrjedEnIq03oC27qjHjympx31un0wiTQLq4bS4shrAECZnMGVTvbD9QSnOhDInjzsyuiBXWQP0&#xA;r4R18jfgnHvJr4nOknrFmJ5xzhraRtHWf3Kkpfz/yGo5S+H73hkquyhdjr3R8m2lj09svwAAAP//&#xA;AwBQSwMEFAAGAAgAAAAhAFr0LFu/AAAAFQEAAAsAAABfcmVscy8ucmVsc2zPwWrDMAwG4Ptg72B0&#xA;X5TuUMaI01uh19I+gLGVxCy2jGSy9e1nemrHjpL4P0nD4SetZiPRyNnCruvBUPYcYp4tXC/Htw8w&#xA;Wl0ObuVMFm6kcBhfX4Yzra62kC6xqGlKVgtLreUTUf1CyWnHhXKbTCzJ1VbKjMX5LzcTvvf9HuXR&#xA;gPHJNKdgQU5hB+ZyK23zHztFL6w81c5zQp6m6P9TMfB3PtPWFCczVQtB9N4U2rp2HOA44NMz4y8A&#xA;AAD//wMAUEsDBBQABgAIAAAAIQBaPFG0wwAAAN0AAAAPAAAAZHJzL2Rvd25yZXYueG1sRI9Bi8JA&#xA;DIXvgv9hiLA3neqyItVRRKh4W7Z68RY7sS12MqUzav335rDgLeG9vPdlteldox7UhdqzgekkAUVc&#xA;eFtzaeB0zMYLUCEiW2w8k4EXBdish4MVptY/+Y8eeSyVhHBI0UAVY5tqHYqKHIaJb4lFu/rOYZS1&#xA;K7Xt8CnhrtGzJJlrhzVLQ4Ut7SoqbvndGbidTz/Z/ndnj02+tZcyi+fL1RrzNeq3S1CR+vgx/18f&#xA;rOB/J8Iv38gIev0GAAD//wMAUEsBAi0AFAAGAAgAAAAhANvh9svuAAAAhQEAABMAAAAAAAAAAAAA&#xA;AAAAAAAAAFtDb250ZW50X1R5cGVzXS54bWxQSwECLQAUAAYACAAAACEAWvQsW78AAAAVAQAACwAA&#xA;AAAAAAAAAAAAAAAfAQAAX3JlbHMvLnJlbHNQSwECLQAUAAYACAAAACEAWjxRtMMAAADdAAAADwAA&#xA;AAAAAAAAAAAAAAAHAgAAZHJzL2Rvd25yZXYueG1sUEsFBgAAAAADAAMAtwAAAPcCAAAAAA==&#xA;" strokeweight="2pt"/><v:line id="Line 928" o:spid="_x0000_s1969" style="position:absolute;visibility:visible" from="9397,14913" to="9400,15450" o:connectortype="straight" o:gfxdata="UEsDBBQABgAIAAAAIQDb4fbL7gAAAIUBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbHyQz07DMAyH&#xA;70i8Q+QralM4IITa7kDhCAiNB7ASt43WOlEcyvb2pNu4IODoPz9/n1xv9vOkForiPDdwXVagiI23&#xA;jocG3rdPxR0oScgWJ8/UwIEENu3lRb09BBKV0ywNjCmFe63FjDSjlD4Q50nv44wpl3HQAc0OB9I3&#xA;VXWrjedEnIq03oC27qjHjympx31un0wiTQLq4bS4shrAECZnMGVTvbD9QSnOhDInjzsyuiBXWQP0&#xA;r4R18jfgnHvJr4nOknrFmJ5xzhraRtHWf3Kkpfz/yGo5S+H73hkquyhdjr3R8m2lj09svwAAAP//&#xA;AwBQSwMEFAAGAAgAAAAhAFr0LFu/AAAAFQEAAAsAAABfcmVscy8ucmVsc2zPwWrDMAwG4Ptg72B0&#xA;X5TuUMaI01uh19I+gLGVxCy2jGSy9e1nemrHjpL4P0nD4SetZiPRyNnCruvBUPYcYp4tXC/Htw8w&#xA;Wl0ObuVMFm6kcBhfX4Yzra62kC6xqGlKVgtLreUTUf1CyWnHhXKbTCzJ1VbKjMX5LzcTvvf9HuXR&#xA;gPHJNKdgQU5hB+ZyK23zHztFL6w81c5zQp6m6P9TMfB3PtPWFCczVQtB9N4U2rp2HOA44NMz4y8A&#xA;AAD//wMAUEsDBBQABgAIAAAAIQC6mWxbvgAAAN0AAAAPAAAAZHJzL2Rvd25yZXYueG1sRE+9CsIw&#xA;EN4F3yGc4KapiiLVKCJU3MTq4nY2Z1tsLqWJWt/eCILbfXy/t1y3phJPalxpWcFoGIEgzqwuOVdw&#xA;PiWDOQjnkTVWlknBmxysV93OEmNtX3ykZ+pzEULYxaig8L6OpXRZQQbd0NbEgbvZxqAPsMmlbvAV&#xA;wk0lx1E0kwZLDg0F1rQtKLunD6PgfjlPk91hq09VutHXPPGX600r1e+1mwUIT63/i3/uvQ7zJ9EM&#xA;vt+EE+TqAwAA//8DAFBLAQItABQABgAIAAAAIQDb4fbL7gAAAIUBAAATAAAAAAAAAAAAAAAAAAAA&#xA;AABbQ29udGVudF9UeXBlc10ueG1sUEsBAi0AFAAGAAgAAAAhAFr0LFu/AAAAFQEAAAsAAAAAAAAA&#xA;AAAAAAAAHwEAAF9yZWxzLy5yZWxzUEsBAi0AFAAGAAgAAAAhALqZbFu+AAAA3QAAAA8AAAAAAAAA&#xA;AAAAAAAABwIAAGRycy9kb3ducmV2LnhtbFBLBQYAAAAAAwADALcAAADyAgAAAAA=&#xA;" strokeweight="2pt"/><v:line id="Line 929" o:spid="_x0000_s1970" style="position:absolute;visibility:visible" from="1129,15730" to="5108,15731" o:connectortype="straight" o:gfxdata="UEsDBBQABgAIAAAAIQDb4fbL7gAAAIUBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbHyQz07DMAyH&#xA;70i8Q+QralM4IITa7kDhCAiNB7ASt43WOlEcyvb2pNu4IODoPz9/n1xv9vOkForiPDdwXVagiI23&#xA;jocG3rdPxR0oScgWJ8/UwIEENu3lRb09BBKV0ywNjCmFe63FjDSjlD4Q50nv44wpl3HQAc0OB9I3&#xA;VXWrjedEnIq03oC27qjHjympx31un0wiTQLq4bS4shrAECZnMGVTvbD9QSnOhDInjzsyuiBXWQP0&#xA;r4R18jfgnHvJr4nOknrFmJ5xzhraRtHWf3Kkpfz/yGo5S+H73hkquyhdjr3R8m2lj09svwAAAP//&#xA;AwBQSwMEFAAGAAgAAAAhAFr0LFu/AAAAFQEAAAsAAABfcmVscy8ucmVsc2zPwWrDMAwG4Ptg72B0&#xA;X5TuUMaI01uh19I+gLGVxCy2jGSy9e1nemrHjpL4P0nD4SetZiPRyNnCruvBUPYcYp4tXC/Htw8w&#xA;Wl0ObuVMFm6kcBhfX4Yzra62kC6xqGlKVgtLreUTUf1CyWnHhXKbTCzJ1VbKjMX5LzcTvvf9HuXR&#xA;gPHJNKdgQU5hB+ZyK23zHztFL6w81c5zQp6m6P9TMfB3PtPWFCczVQtB9N4U2rp2HOA44NMz4y8A&#xA;AAD//wMAUEsDBBQABgAIAAAAIQBVIKGDxgAAAN0AAAAPAAAAZHJzL2Rvd25yZXYueG1sRI9BTwIx&#xA;EIXvJv6HZky4SXchMbpQCFFMJB6MwA8YtsN2YTvdtBVWf71zMPE2k/fmvW/my8F36kIxtYENlOMC&#xA;FHEdbMuNgf3u9f4RVMrIFrvAZOCbEiwXtzdzrGy48iddtrlREsKpQgMu577SOtWOPKZx6IlFO4bo&#xA;McsaG20jXiXcd3pSFA/aY8vS4LCnZ0f1efvlDWzi4f1c/jROH3gT193Hy1PyJ2NGd8NqBirTkP/N&#xA;f9dvVvCnpfDLNzKCXvwCAAD//wMAUEsBAi0AFAAGAAgAAAAhANvh9svuAAAAhQEAABMAAAAAAAAA&#xA;AAAAAAAAAAAAAFtDb250ZW50X1R5cGVzXS54bWxQSwECLQAUAAYACAAAACEAWvQsW78AAAAVAQAA&#xA;CwAAAAAAAAAAAAAAAAAfAQAAX3JlbHMvLnJlbHNQSwECLQAUAAYACAAAACEAVSChg8YAAADdAAAA&#xA;DwAAAAAAAAAAAAAAAAAHAgAAZHJzL2Rvd25yZXYueG1sUEsFBgAAAAADAAMAtwAAAPoCAAAAAA==&#xA;" strokeweight="1pt"/><v:line id="Line 930" o:spid="_x0000_s1971" style="position:absolute;visibility:visible" from="1129,16004" to="5108,16004" o:connectortype="straight" o:gfxdata="UEsDBBQABgAIAAAAIQDb4fbL7gAAAIUBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbHyQz07DMAyH&#xA;70i8Q+QralM4IITa7kDhCAiNB7ASt43WOlEcyvb2pNu4IODoPz9/n1xv9vOkForiPDdwXVagiI23&#xA;jocG3rdPxR0oScgWJ8/UwIEENu3lRb09BBKV0ywNjCmFe63FjDSjlD4Q50nv44wpl3HQAc0OB9I3&#xA;VXWrjedEnIq03oC27qjHjympx31un0wiTQLq4bS4shrAECZnMGVTvbD9QSnOhDInjzsyuiBXWQP0&#xA;r4R18jfgnHvJr4nOknrFmJ5xzhraRtHWf3Kkpfz/yGo5S+H73hkquyhdjr3R8m2lj09svwAAAP//&#xA;AwBQSwMEFAAGAAgAAAAhAFr0LFu/AAAAFQEAAAsAAABfcmVscy8ucmVsc2zPwWrDMAwG4Ptg72B0&#xA;X5TuUMaI01uh19I+gLGVxCy2jGSy9e1nemrHjpL4P0nD4SetZiPRyNnCruvBUPYcYp4tXC/Htw8w&#xA;Wl0ObuVMFm6kcBhfX4Yzra62kC6xqGlKVgtLreUTUf1CyWnHhXKbTCzJ1VbKjMX5LzcTvvf9HuXR&#xA;gPHJNKdgQU5hB+ZyK23zHztFL6w81c5zQp6m6P9TMfB3PtPWFCczVQtB9N4U2rp2HOA44NMz4y8A&#xA;AAD//wMAUEsDBBQABgAIAAAAIQA6bAQYwwAAAN0AAAAPAAAAZHJzL2Rvd25yZXYueG1sRE/NagIx&#xA;EL4XfIcwgjfNrgWpW6OItqB4KGofYNxMN1s3kyVJde3TN4LQ23x8vzNbdLYRF/KhdqwgH2UgiEun&#xA;a64UfB7fhy8gQkTW2DgmBTcKsJj3nmZYaHflPV0OsRIphEOBCkyMbSFlKA1ZDCPXEifuy3mLMUFf&#xA;Se3xmsJtI8dZNpEWa04NBltaGSrPhx+rYOtPu3P+Wxl54q1/az7W02C/lRr0u+UriEhd/Bc/3Bud&#xA;5j/nOdy/SSfI+R8AAAD//wMAUEsBAi0AFAAGAAgAAAAhANvh9svuAAAAhQEAABMAAAAAAAAAAAAA&#xA;AAAAAAAAAFtDb250ZW50X1R5cGVzXS54bWxQSwECLQAUAAYACAAAACEAWvQsW78AAAAVAQAACwAA&#xA;AAAAAAAAAAAAAAAfAQAAX3JlbHMvLnJlbHNQSwECLQAUAAYACAAAACEAOmwEGMMAAADdAAAADwAA&#xA;AAAAAAAAAAAAAAAHAgAAZHJzL2Rvd25yZXYueG1sUEsFBgAAAAADAAMAtwAAAPcCAAAAAA==&#xA;" strokeweight="1pt"/><v:rect id="Rectangle 931" o:spid="_x0000_s1972" style="position:absolute;left:1152;top:14659;width:461;height:240;visibility:visible" o:gfxdata="UEsDBBQABgAIAAAAIQDb4fbL7gAAAIUBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbHyQz07DMAyH&#xA;70i8Q+QralM4IITa7kDhCAiNB7ASt43WOlEcyvb2pNu4IODoPz9/n1xv9vOkForiPDdwXVagiI23&#xA;jocG3rdPxR0oScgWJ8/UwIEENu3lRb09BBKV0ywNjCmFe63FjDSjlD4Q50nv44wpl3HQAc0OB9I3&#xA;VXWrjedEnIq03oC27qjHjympx31un0wiTQLq4bS4shrAECZnMGVTvbD9QSnOhDInjzsyuiBXWQP0&#xA;r4R18jfgnHvJr4nOknrFmJ5xzhraRtHWf3Kkpfz/yGo5S+H73hkquyhdjr3R8m2lj09svwAAAP//&#xA;AwBQSwMEFAAGAAgAAAAhAFr0LFu/AAAAFQEAAAsAAABfcmVscy8ucmVsc2zPwWrDMAwG4Ptg72B0&#xA;X5TuUMaI01uh19I+gLGVxCy2jGSy9e1nemrHjpL4P0nD4SetZiPRyNnCruvBUPYcYp4tXC/Htw8w&#xA;Wl0ObuVMFm6kcBhfX4Yzra62kC6xqGlKVgtLreUTUf1CyWnHhXKbTCzJ1VbKjMX5LzcTvvf9HuXR&#xA;gPHJNKdgQU5hB+ZyK23zHztFL6w81c5zQp6m6P9TMfB3PtPWFCczVQtB9N4U2rp2HOA44NMz4y8A&#xA;AAD//wMAUEsDBBQABgAIAAAAIQDUID+qwQAAAN0AAAAPAAAAZHJzL2Rvd25yZXYueG1sRE/Pa8Iw&#xA;FL4P/B/CE7ytqUU61zVKEYRd7SZ4fDRvbWfzUpPM1v9+OQx2/Ph+l/vZDOJOzveWFayTFARxY3XP&#xA;rYLPj+PzFoQPyBoHy6TgQR72u8VTiYW2E5/oXodWxBD2BSroQhgLKX3TkUGf2JE4cl/WGQwRulZq&#xA;h1MMN4PM0jSXBnuODR2OdOioudY/RkFVfc/nW/2KRy+3qcv1RrfVRanVcq7eQASaw7/4z/2uFWQv&#xA;Wdwf38QnIHe/AAAA//8DAFBLAQItABQABgAIAAAAIQDb4fbL7gAAAIUBAAATAAAAAAAAAAAAAAAA&#xA;AAAAAABbQ29udGVudF9UeXBlc10ueG1sUEsBAi0AFAAGAAgAAAAhAFr0LFu/AAAAFQEAAAsAAAAA&#xA;AAAAAAAAAAAAHwEAAF9yZWxzLy5yZWxzUEsBAi0AFAAGAAgAAAAhANQgP6rBAAAA3QAAAA8AAAAA&#xA;AAAAAAAAAAAABwIAAGRycy9kb3ducmV2LnhtbFBLBQYAAAAAAwADALcAAAD1AgAAAAA=&#xA;" filled="f" stroked="f" strokeweight=".25pt"><v:textbox style="mso-next-textbox:#Rectangle 931" inset="1pt,1pt,1pt,1pt"><w:txbxContent><w:p w14:paraId="7A042E10" w14:textId="77777777" w:rsidR="00B827FD" w:rsidRPr="00951A0B" w:rsidRDefault="00B827FD" w:rsidP="00B827FD"><w:pPr><w:pStyle w:val="a1"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr><w:r w:rsidRPr="00951A0B"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>Изм.</w:t></w:r></w:p></w:txbxContent></v:textbox></v:rect><v:rect id="Rectangle 932" o:spid="_x0000_s1973" style="position:absolute;left:1673;top:14659;width:574;height:240;visibility:visible" o:gfxdata="UEsDBBQABgAIAAAAIQDb4fbL7gAAAIUBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbHyQz07DMAyH&#xA;70i8Q+QralM4IITa7kDhCAiNB7ASt43WOlEcyvb2pNu4IODoPz9/n1xv9vOkForiPDdwXVagiI23&#xA;jocG3rdPxR0oScgWJ8/UwIEENu3lRb09BBKV0ywNjCmFe63FjDSjlD4Q50nv44wpl3HQAc0OB9I3&#xA;VXWrjedEnIq03oC27qjHjympx31un0wiTQLq4bS4shrAECZnMGVTvbD9QSnOhDInjzsyuiBXWQP0&#xA;r4R18jfgnHvJr4nOknrFmJ5xzhraRtHWf3Kkpfz/yGo5S+H73hkquyhdjr3R8m2lj09svwAAAP//&#xA;AwBQSwMEFAAGAAgAAAAhAFr0LFu/AAAAFQEAAAsAAABfcmVscy8ucmVsc2zPwWrDMAwG4Ptg72B0&#xA;X5TuUMaI01uh19I+gLGVxCy2jGSy9e1nemrHjpL4P0nD4SetZiPRyNnCruvBUPYcYp4tXC/Htw8w&#xA;Wl0ObuVMFm6kcBhfX4Yzra62kC6xqGlKVgtLreUTUf1CyWnHhXKbTCzJ1VbKjMX5LzcTvvf9HuXR&#xA;gPHJNKdgQU5hB+ZyK23zHztFL6w81c5zQp6m6P9TMfB3PtPWFCczVQtB9N4U2rp2HOA44NMz4y8A&#xA;AAD//wMAUEsDBBQABgAIAAAAIQC7bJoxxAAAAN0AAAAPAAAAZHJzL2Rvd25yZXYueG1sRI/BasMw&#xA;EETvhfyD2EButRxTXMeJEkwh0GvdFnpcrI3txFo5kho7f18VCj0OM/OG2R1mM4gbOd9bVrBOUhDE&#xA;jdU9two+3o+PBQgfkDUOlknBnTwc9ouHHZbaTvxGtzq0IkLYl6igC2EspfRNRwZ9Ykfi6J2sMxii&#xA;dK3UDqcIN4PM0jSXBnuOCx2O9NJRc6m/jYKqOs+f13qDRy+L1OX6SbfVl1Kr5VxtQQSaw3/4r/2q&#xA;FWTP2Rp+38QnIPc/AAAA//8DAFBLAQItABQABgAIAAAAIQDb4fbL7gAAAIUBAAATAAAAAAAAAAAA&#xA;AAAAAAAAAABbQ29udGVudF9UeXBlc10ueG1sUEsBAi0AFAAGAAgAAAAhAFr0LFu/AAAAFQEAAAsA&#xA;AAAAAAAAAAAAAAAAHwEAAF9yZWxzLy5yZWxzUEsBAi0AFAAGAAgAAAAhALtsmjHEAAAA3QAAAA8A&#xA;AAAAAAAAAAAAAAAABwIAAGRycy9kb3ducmV2LnhtbFBLBQYAAAAAAwADALcAAAD4AgAAAAA=&#xA;" filled="f" stroked="f" strokeweight=".25pt"><v:textbox style="mso-next-textbox:#Rectangle 932" inset="1pt,1pt,1pt,1pt"><w:txbxContent><w:p w14:paraId="0D680384" w14:textId="77777777" w:rsidR="00B827FD" w:rsidRPr="00017D2F" w:rsidRDefault="00B827FD" w:rsidP="00B827FD"><w:pPr><w:pStyle w:val="a1"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="ru-RU"/></w:rPr></w:pPr><w:r w:rsidRPr="00017D2F"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="ru-RU"/></w:rPr><w:t>Лист</w:t></w:r></w:p></w:txbxContent></v:textbox></v:rect><v:rect id="Rectangle 933" o:spid="_x0000_s1974" style="position:absolute;left:2307;top:14659;width:1344;height:240;visibility:visible" o:gfxdata="UEsDBBQABgAIAAAAIQDb4fbL7gAAAIUBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbHyQz07DMAyH&#xA;70i8Q+QralM4IITa7kDhCAiNB7ASt43WOlEcyvb2pNu4IODoPz9/n1xv9vOkForiPDdwXVagiI23&#xA;jocG3rdPxR0oScgWJ8/UwIEENu3lRb09BBKV0ywNjCmFe63FjDSjlD4Q50nv44wpl3HQAc0OB9I3&#xA;VXWrjedEnIq03oC27qjHjympx31un0wiTQLq4bS4shrAECZnMGVTvbD9QSnOhDInjzsyuiBXWQP0&#xA;r4R18jfgnHvJr4nOknrFmJ5xzhraRtHWf3Kkpfz/yGo5S+H73hkquyhdjr3R8m2lj09svwAAAP//&#xA;AwBQSwMEFAAGAAgAAAAhAFr0LFu/AAAAFQEAAAsAAABfcmVscy8ucmVsc2zPwWrDMAwG4Ptg72B0&#xA;X5TuUMaI01uh19I+gLGVxCy2jGSy9e1nemrHjpL4P0nD4SetZiPRyNnCruvBUPYcYp4tXC/Htw8w&#xA;Wl0ObuVMFm6kcBhfX4Yzra62kC6xqGlKVgtLreUTUf1CyWnHhXKbTCzJ1VbKjMX5LzcTvvf9HuXR&#xA;gPHJNKdgQU5hB+ZyK23zHztFL6w81c5zQp6m6P9TMfB3PtPWFCczVQtB9N4U2rp2HOA44NMz4y8A&#xA;AAD//wMAUEsDBBQABgAIAAAAIQBLvgRGwwAAAN0AAAAPAAAAZHJzL2Rvd25yZXYueG1sRI9Pi8Iw&#xA;FMTvwn6H8Bb2pukW8U/XKEUQvFoVPD6at221eekmUbvf3giCx2FmfsMsVr1pxY2cbywr+B4lIIhL&#xA;qxuuFBz2m+EMhA/IGlvLpOCfPKyWH4MFZtreeUe3IlQiQthnqKAOocuk9GVNBv3IdsTR+7XOYIjS&#xA;VVI7vEe4aWWaJBNpsOG4UGNH65rKS3E1CvL83B//ijluvJwlbqLHuspPSn199vkPiEB9eIdf7a1W&#xA;kE7TFJ5v4hOQywcAAAD//wMAUEsBAi0AFAAGAAgAAAAhANvh9svuAAAAhQEAABMAAAAAAAAAAAAA&#xA;AAAAAAAAAFtDb250ZW50X1R5cGVzXS54bWxQSwECLQAUAAYACAAAACEAWvQsW78AAAAVAQAACwAA&#xA;AAAAAAAAAAAAAAAfAQAAX3JlbHMvLnJlbHNQSwECLQAUAAYACAAAACEAS74ERsMAAADdAAAADwAA&#xA;AAAAAAAAAAAAAAAHAgAAZHJzL2Rvd25yZXYueG1sUEsFBgAAAAADAAMAtwAAAPcCAAAAAA==&#xA;" filled="f" stroked="f" strokeweight=".25pt"><v:textbox style="mso-next-textbox:#Rectangle 933" inset="1pt,1pt,1pt,1pt"><w:txbxContent><w:p w14:paraId="6A5BA16A" w14:textId="77777777" w:rsidR="00B827FD" w:rsidRPr="00017D2F" w:rsidRDefault="00B827FD" w:rsidP="00B827FD"><w:pPr><w:pStyle w:val="a1"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr><w:r w:rsidRPr="00017D2F"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>№ докум.</w:t></w:r></w:p></w:txbxContent></v:textbox></v:rect><v:rect id="Rectangle 934" o:spid="_x0000_s1975" style="position:absolute;left:3725;top:14659;width:801;height:240;visibility:visible" o:gfxdata="UEsDBBQABgAIAAAAIQDb4fbL7gAAAIUBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbHyQz07DMAyH&#xA;70i8Q+QralM4IITa7kDhCAiNB7ASt43WOlEcyvb2pNu4IODoPz9/n1xv9vOkForiPDdwXVagiI23&#xA;jocG3rdPxR0oScgWJ8/UwIEENu3lRb09BBKV0ywNjCmFe63FjDSjlD4Q50nv44wpl3HQAc0OB9I3&#xA;VXWrjedEnIq03oC27qjHjympx31un0wiTQLq4bS4shrAECZnMGVTvbD9QSnOhDInjzsyuiBXWQP0&#xA;r4R18jfgnHvJr4nOknrFmJ5xzhraRtHWf3Kkpfz/yGo5S+H73hkquyhdjr3R8m2lj09svwAAAP//&#xA;AwBQSwMEFAAGAAgAAAAhAFr0LFu/AAAAFQEAAAsAAABfcmVscy8ucmVsc2zPwWrDMAwG4Ptg72B0&#xA;X5TuUMaI01uh19I+gLGVxCy2jGSy9e1nemrHjpL4P0nD4SetZiPRyNnCruvBUPYcYp4tXC/Htw8w&#xA;Wl0ObuVMFm6kcBhfX4Yzra62kC6xqGlKVgtLreUTUf1CyWnHhXKbTCzJ1VbKjMX5LzcTvvf9HuXR&#xA;gPHJNKdgQU5hB+ZyK23zHztFL6w81c5zQp6m6P9TMfB3PtPWFCczVQtB9N4U2rp2HOA44NMz4y8A&#xA;AAD//wMAUEsDBBQABgAIAAAAIQAk8qHdwwAAAN0AAAAPAAAAZHJzL2Rvd25yZXYueG1sRI9Ba8JA&#xA;FITvhf6H5RV6qxuj2DS6ShCEXk0VenxkX5No9m3cXTX+e1cQehxm5htmsRpMJy7kfGtZwXiUgCCu&#xA;rG65VrD72XxkIHxA1thZJgU38rBavr4sMNf2ylu6lKEWEcI+RwVNCH0upa8aMuhHtieO3p91BkOU&#xA;rpba4TXCTSfTJJlJgy3HhQZ7WjdUHcuzUVAUh2F/Kr9w42WWuJme6rr4Ver9bSjmIAIN4T/8bH9r&#xA;BelnOoHHm/gE5PIOAAD//wMAUEsBAi0AFAAGAAgAAAAhANvh9svuAAAAhQEAABMAAAAAAAAAAAAA&#xA;AAAAAAAAAFtDb250ZW50X1R5cGVzXS54bWxQSwECLQAUAAYACAAAACEAWvQsW78AAAAVAQAACwAA&#xA;AAAAAAAAAAAAAAAfAQAAX3JlbHMvLnJlbHNQSwECLQAUAAYACAAAACEAJPKh3cMAAADdAAAADwAA&#xA;AAAAAAAAAAAAAAAHAgAAZHJzL2Rvd25yZXYueG1sUEsFBgAAAAADAAMAtwAAAPcCAAAAAA==&#xA;" filled="f" stroked="f" strokeweight=".25pt"><v:textbox style="mso-next-textbox:#Rectangle 934" inset="1pt,1pt,1pt,1pt"><w:txbxContent><w:p w14:paraId="584C3DC4" w14:textId="77777777" w:rsidR="00B827FD" w:rsidRPr="00017D2F" w:rsidRDefault="00B827FD" w:rsidP="00B827FD"><w:pPr><w:pStyle w:val="a1"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr><w:r w:rsidRPr="00017D2F"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>Подпись</w:t></w:r></w:p></w:txbxContent></v:textbox></v:rect><v:rect id="Rectangle 935" o:spid="_x0000_s1976" style="position:absolute;left:4572;top:14659;width:522;height:240;visibility:visible" o:gfxdata="UEsDBBQABgAIAAAAIQDb4fbL7gAAAIUBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbHyQz07DMAyH&#xA;70i8Q+QralM4IITa7kDhCAiNB7ASt43WOlEcyvb2pNu4IODoPz9/n1xv9vOkForiPDdwXVagiI23&#xA;jocG3rdPxR0oScgWJ8/UwIEENu3lRb09BBKV0ywNjCmFe63FjDSjlD4Q50nv44wpl3HQAc0OB9I3&#xA;VXWrjedEnIq03oC27qjHjympx31un0wiTQLq4bS4shrAECZnMGVTvbD9QSnOhDInjzsyuiBXWQP0&#xA;r4R18jfgnHvJr4nOknrFmJ5xzhraRtHWf3Kkpfz/yGo5S+H73hkquyhdjr3R8m2lj09svwAAAP//&#xA;AwBQSwMEFAAGAAgAAAAhAFr0LFu/AAAAFQEAAAsAAABfcmVscy8ucmVsc2zPwWrDMAwG4Ptg72B0&#xA;X5TuUMaI01uh19I+gLGVxCy2jGSy9e1nemrHjpL4P0nD4SetZiPRyNnCruvBUPYcYp4tXC/Htw8w&#xA;Wl0ObuVMFm6kcBhfX4Yzra62kC6xqGlKVgtLreUTUf1CyWnHhXKbTCzJ1VbKjMX5LzcTvvf9HuXR&#xA;gPHJNKdgQU5hB+ZyK23zHztFL6w81c5zQp6m6P9TMfB3PtPWFCczVQtB9N4U2rp2HOA44NMz4y8A&#xA;AAD//wMAUEsDBBQABgAIAAAAIQCrGzmpxAAAAN0AAAAPAAAAZHJzL2Rvd25yZXYueG1sRI/BasMw&#xA;EETvhfyD2EButVxjXNeJEkwhkGvdBnJcrK3t1Fq5kpI4f18VCj0OM/OG2exmM4orOT9YVvCUpCCI&#xA;W6sH7hR8vO8fSxA+IGscLZOCO3nYbRcPG6y0vfEbXZvQiQhhX6GCPoSpktK3PRn0iZ2Io/dpncEQ&#xA;peukdniLcDPKLE0LaXDguNDjRK89tV/NxSio6/N8/G5ecO9lmbpC57qrT0qtlnO9BhFoDv/hv/ZB&#xA;K8iesxx+38QnILc/AAAA//8DAFBLAQItABQABgAIAAAAIQDb4fbL7gAAAIUBAAATAAAAAAAAAAAA&#xA;AAAAAAAAAABbQ29udGVudF9UeXBlc10ueG1sUEsBAi0AFAAGAAgAAAAhAFr0LFu/AAAAFQEAAAsA&#xA;AAAAAAAAAAAAAAAAHwEAAF9yZWxzLy5yZWxzUEsBAi0AFAAGAAgAAAAhAKsbOanEAAAA3QAAAA8A&#xA;AAAAAAAAAAAAAAAABwIAAGRycy9kb3ducmV2LnhtbFBLBQYAAAAAAwADALcAAAD4AgAAAAA=&#xA;" filled="f" stroked="f" strokeweight=".25pt"><v:textbox style="mso-next-textbox:#Rectangle 935" inset="1pt,1pt,1pt,1pt"><w:txbxContent><w:p w14:paraId="7FEA755F" w14:textId="77777777" w:rsidR="00B827FD" w:rsidRPr="00017D2F" w:rsidRDefault="00B827FD" w:rsidP="00B827FD"><w:pPr><w:pStyle w:val="a1"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:spacing w:val="-10"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr><w:r w:rsidRPr="00017D2F"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:spacing w:val="-10"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>Дата</w:t></w:r></w:p></w:txbxContent></v:textbox></v:rect><v:rect id="Rectangle 936" o:spid="_x0000_s1977" style="position:absolute;left:9440;top:14928;width:770;height:239;visibility:visible" o:gfxdata="UEsDBBQABgAIAAAAIQDb4fbL7gAAAIUBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbHyQz07DMAyH&#xA;70i8Q+QralM4IITa7kDhCAiNB7ASt43WOlEcyvb2pNu4IODoPz9/n1xv9vOkForiPDdwXVagiI23&#xA;jocG3rdPxR0oScgWJ8/UwIEENu3lRb09BBKV0ywNjCmFe63FjDSjlD4Q50nv44wpl3HQAc0OB9I3&#xA;VXWrjedEnIq03oC27qjHjympx31un0wiTQLq4bS4shrAECZnMGVTvbD9QSnOhDInjzsyuiBXWQP0&#xA;r4R18jfgnHvJr4nOknrFmJ5xzhraRtHWf3Kkpfz/yGo5S+H73hkquyhdjr3R8m2lj09svwAAAP//&#xA;AwBQSwMEFAAGAAgAAAAhAFr0LFu/AAAAFQEAAAsAAABfcmVscy8ucmVsc2zPwWrDMAwG4Ptg72B0&#xA;X5TuUMaI01uh19I+gLGVxCy2jGSy9e1nemrHjpL4P0nD4SetZiPRyNnCruvBUPYcYp4tXC/Htw8w&#xA;Wl0ObuVMFm6kcBhfX4Yzra62kC6xqGlKVgtLreUTUf1CyWnHhXKbTCzJ1VbKjMX5LzcTvvf9HuXR&#xA;gPHJNKdgQU5hB+ZyK23zHztFL6w81c5zQp6m6P9TMfB3PtPWFCczVQtB9N4U2rp2HOA44NMz4y8A&#xA;AAD//wMAUEsDBBQABgAIAAAAIQDEV5wywwAAAN0AAAAPAAAAZHJzL2Rvd25yZXYueG1sRI9Ba8JA&#xA;FITvhf6H5RV6qxuD2jS6ShCEXk0VenxkX5No9m3cXTX+e1cQehxm5htmsRpMJy7kfGtZwXiUgCCu&#xA;rG65VrD72XxkIHxA1thZJgU38rBavr4sMNf2ylu6lKEWEcI+RwVNCH0upa8aMuhHtieO3p91BkOU&#xA;rpba4TXCTSfTJJlJgy3HhQZ7WjdUHcuzUVAUh2F/Kr9w42WWuJme6Lr4Ver9bSjmIAIN4T/8bH9r&#xA;BelnOoXHm/gE5PIOAAD//wMAUEsBAi0AFAAGAAgAAAAhANvh9svuAAAAhQEAABMAAAAAAAAAAAAA&#xA;AAAAAAAAAFtDb250ZW50X1R5cGVzXS54bWxQSwECLQAUAAYACAAAACEAWvQsW78AAAAVAQAACwAA&#xA;AAAAAAAAAAAAAAAfAQAAX3JlbHMvLnJlbHNQSwECLQAUAAYACAAAACEAxFecMsMAAADdAAAADwAA&#xA;AAAAAAAAAAAAAAAHAgAAZHJzL2Rvd25yZXYueG1sUEsFBgAAAAADAAMAtwAAAPcCAAAAAA==&#xA;" filled="f" stroked="f" strokeweight=".25pt"><v:textbox style="mso-next-textbox:#Rectangle 936" inset="1pt,1pt,1pt,1pt"><w:txbxContent><w:p w14:paraId="60A5B04D" w14:textId="77777777" w:rsidR="00B827FD" w:rsidRPr="00B827FD" w:rsidRDefault="00B827FD" w:rsidP="00B827FD"><w:pPr><w:pStyle w:val="a1"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr><w:r w:rsidRPr="00B827FD"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>Лист</w:t></w:r></w:p></w:txbxContent></v:textbox></v:rect><v:rect id="Rectangle 937" o:spid="_x0000_s1978" style="position:absolute;left:9440;top:15211;width:770;height:240;visibility:visible" o:gfxdata="UEsDBBQABgAIAAAAIQDb4fbL7gAAAIUBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbHyQz07DMAyH&#xA;70i8Q+QralM4IITa7kDhCAiNB7ASt43WOlEcyvb2pNu4IODoPz9/n1xv9vOkForiPDdwXVagiI23&#xA;jocG3rdPxR0oScgWJ8/UwIEENu3lRb09BBKV0ywNjCmFe63FjDSjlD4Q50nv44wpl3HQAc0OB9I3&#xA;VXWrjedEnIq03oC27qjHjympx31un0wiTQLq4bS4shrAECZnMGVTvbD9QSnOhDInjzsyuiBXWQP0&#xA;r4R18jfgnHvJr4nOknrFmJ5xzhraRtHWf3Kkpfz/yGo5S+H73hkquyhdjr3R8m2lj09svwAAAP//&#xA;AwBQSwMEFAAGAAgAAAAhAFr0LFu/AAAAFQEAAAsAAABfcmVscy8ucmVsc2zPwWrDMAwG4Ptg72B0&#xA;X5TuUMaI01uh19I+gLGVxCy2jGSy9e1nemrHjpL4P0nD4SetZiPRyNnCruvBUPYcYp4tXC/Htw8w&#xA;Wl0ObuVMFm6kcBhfX4Yzra62kC6xqGlKVgtLreUTUf1CyWnHhXKbTCzJ1VbKjMX5LzcTvvf9HuXR&#xA;gPHJNKdgQU5hB+ZyK23zHztFL6w81c5zQp6m6P9TMfB3PtPWFCczVQtB9N4U2rp2HOA44NMz4y8A&#xA;AAD//wMAUEsDBBQABgAIAAAAIQA0hQJFwwAAAN0AAAAPAAAAZHJzL2Rvd25yZXYueG1sRI9Ba8JA&#xA;FITvhf6H5RW81U1DSTW6CaEg9Gqs4PGRfSbR7Nt0d6vx37uFgsdhZr5h1uVkBnEh53vLCt7mCQji&#xA;xuqeWwXfu83rAoQPyBoHy6TgRh7K4vlpjbm2V97SpQ6tiBD2OSroQhhzKX3TkUE/tyNx9I7WGQxR&#xA;ulZqh9cIN4NMkySTBnuOCx2O9NlRc65/jYKqOk37n3qJGy8Xicv0u26rg1Kzl6lagQg0hUf4v/2l&#xA;FaQfaQZ/b+ITkMUdAAD//wMAUEsBAi0AFAAGAAgAAAAhANvh9svuAAAAhQEAABMAAAAAAAAAAAAA&#xA;AAAAAAAAAFtDb250ZW50X1R5cGVzXS54bWxQSwECLQAUAAYACAAAACEAWvQsW78AAAAVAQAACwAA&#xA;AAAAAAAAAAAAAAAfAQAAX3JlbHMvLnJlbHNQSwECLQAUAAYACAAAACEANIUCRcMAAADdAAAADwAA&#xA;AAAAAAAAAAAAAAAHAgAAZHJzL2Rvd25yZXYueG1sUEsFBgAAAAADAAMAtwAAAPcCAAAAAA==&#xA;" filled="f" stroked="f" strokeweight=".25pt"><v:textbox style="mso-next-textbox:#Rectangle 937" inset="1pt,1pt,1pt,1pt"><w:txbxContent><w:p w14:paraId="4F695E4D" w14:textId="77777777" w:rsidR="00B827FD" w:rsidRPr="00B827FD" w:rsidRDefault="00B827FD" w:rsidP="00B827FD"><w:pPr><w:pStyle w:val="a1"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="ru-RU"/></w:rPr></w:pPr><w:r w:rsidRPr="00B827FD"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="ru-RU"/></w:rPr><w:t>1</w:t></w:r></w:p></w:txbxContent></v:textbox></v:rect><v:rect id="Rectangle 938" o:spid="_x0000_s1979" style="position:absolute;left:5175;top:14245;width:6348;height:512;visibility:visible" o:gfxdata="UEsDBBQABgAIAAAAIQDb4fbL7gAAAIUBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbHyQz07DMAyH&#xA;70i8Q+QralM4IITa7kDhCAiNB7ASt43WOlEcyvb2pNu4IODoPz9/n1xv9vOkForiPDdwXVagiI23&#xA;jocG3rdPxR0oScgWJ8/UwIEENu3lRb09BBKV0ywNjCmFe63FjDSjlD4Q50nv44wpl3HQAc0OB9I3&#xA;VXWrjedEnIq03oC27qjHjympx31un0wiTQLq4bS4shrAECZnMGVTvbD9QSnOhDInjzsyuiBXWQP0&#xA;r4R18jfgnHvJr4nOknrFmJ5xzhraRtHWf3Kkpfz/yGo5S+H73hkquyhdjr3R8m2lj09svwAAAP//&#xA;AwBQSwMEFAAGAAgAAAAhAFr0LFu/AAAAFQEAAAsAAABfcmVscy8ucmVsc2zPwWrDMAwG4Ptg72B0&#xA;X5TuUMaI01uh19I+gLGVxCy2jGSy9e1nemrHjpL4P0nD4SetZiPRyNnCruvBUPYcYp4tXC/Htw8w&#xA;Wl0ObuVMFm6kcBhfX4Yzra62kC6xqGlKVgtLreUTUf1CyWnHhXKbTCzJ1VbKjMX5LzcTvvf9HuXR&#xA;gPHJNKdgQU5hB+ZyK23zHztFL6w81c5zQp6m6P9TMfB3PtPWFCczVQtB9N4U2rp2HOA44NMz4y8A&#xA;AAD//wMAUEsDBBQABgAIAAAAIQBbyafexAAAAN0AAAAPAAAAZHJzL2Rvd25yZXYueG1sRI/BasMw&#xA;EETvhfyD2EBvtRxTbMeJEkwh0GvcFnpcrI3txFq5kpo4f18VCj0OM/OG2e5nM4orOT9YVrBKUhDE&#xA;rdUDdwre3w5PJQgfkDWOlknBnTzsd4uHLVba3vhI1yZ0IkLYV6igD2GqpPRtTwZ9Yifi6J2sMxii&#xA;dJ3UDm8RbkaZpWkuDQ4cF3qc6KWn9tJ8GwV1fZ4/vpo1HrwsU5frZ93Vn0o9Lud6AyLQHP7Df+1X&#xA;rSArsgJ+38QnIHc/AAAA//8DAFBLAQItABQABgAIAAAAIQDb4fbL7gAAAIUBAAATAAAAAAAAAAAA&#xA;AAAAAAAAAABbQ29udGVudF9UeXBlc10ueG1sUEsBAi0AFAAGAAgAAAAhAFr0LFu/AAAAFQEAAAsA&#xA;AAAAAAAAAAAAAAAAHwEAAF9yZWxzLy5yZWxzUEsBAi0AFAAGAAgAAAAhAFvJp97EAAAA3QAAAA8A&#xA;AAAAAAAAAAAAAAAABwIAAGRycy9kb3ducmV2LnhtbFBLBQYAAAAAAwADALcAAAD4AgAAAAA=&#xA;" filled="f" stroked="f" strokeweight=".25pt"><v:textbox style="mso-next-textbox:#Rectangle 938" inset="1pt,1pt,1pt,1pt"><w:txbxContent><w:p w14:paraId="18539222" w14:textId="77777777" w:rsidR="00B827FD" w:rsidRPr="008F5AF2" w:rsidRDefault="00B827FD" w:rsidP="00B827FD"><w:pPr><w:pStyle w:val="a1"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="28"/><w:lang w:val="ru-RU"/></w:rPr></w:pPr><w:r w:rsidRPr="0083149A"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Б</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>ГТУ</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="ru-RU"/></w:rPr><w:t xml:space="preserve"> ДП</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> 05.01.ГЧ</w:t></w:r></w:p><w:p w14:paraId="5919D3E1" w14:textId="77777777" w:rsidR="00B827FD" w:rsidRPr="0042148F" w:rsidRDefault="00B827FD" w:rsidP="00B827FD"><w:pPr><w:pStyle w:val="a1"/><w:jc w:val="left"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:szCs w:val="16"/><w:lang w:val="ru-RU"/></w:rPr></w:pPr></w:p><w:p w14:paraId="4EE23B0D" w14:textId="77777777" w:rsidR="00B827FD" w:rsidRPr="0042148F" w:rsidRDefault="00B827FD" w:rsidP="00B827FD"><w:pPr><w:rPr><w:i/><w:szCs w:val="16"/><w:lang w:val="ru-RU"/></w:rPr></w:pPr></w:p></w:txbxContent></v:textbox></v:rect><v:line id="Line 939" o:spid="_x0000_s1980" style="position:absolute;visibility:visible" from="1130,14908" to="11554,14909" o:connectortype="straight" o:gfxdata="UEsDBBQABgAIAAAAIQDb4fbL7gAAAIUBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbHyQz07DMAyH&#xA;70i8Q+QralM4IITa7kDhCAiNB7ASt43WOlEcyvb2pNu4IODoPz9/n1xv9vOkForiPDdwXVagiI23&#xA;jocG3rdPxR0oScgWJ8/UwIEENu3lRb09BBKV0ywNjCmFe63FjDSjlD4Q50nv44wpl3HQAc0OB9I3&#xA;VXWrjedEnIq03oC27qjHjympx31un0wiTQLq4bS4shrAECZnMGVTvbD9QSnOhDInjzsyuiBXWQP0&#xA;r4R18jfgnHvJr4nOknrFmJ5xzhraRtHWf3Kkpfz/yGo5S+H73hkquyhdjr3R8m2lj09svwAAAP//&#xA;AwBQSwMEFAAGAAgAAAAhAFr0LFu/AAAAFQEAAAsAAABfcmVscy8ucmVsc2zPwWrDMAwG4Ptg72B0&#xA;X5TuUMaI01uh19I+gLGVxCy2jGSy9e1nemrHjpL4P0nD4SetZiPRyNnCruvBUPYcYp4tXC/Htw8w&#xA;Wl0ObuVMFm6kcBhfX4Yzra62kC6xqGlKVgtLreUTUf1CyWnHhXKbTCzJ1VbKjMX5LzcTvvf9HuXR&#xA;gPHJNKdgQU5hB+ZyK23zHztFL6w81c5zQp6m6P9TMfB3PtPWFCczVQtB9N4U2rp2HOA44NMz4y8A&#xA;AAD//wMAUEsDBBQABgAIAAAAIQCNGtDavgAAAN0AAAAPAAAAZHJzL2Rvd25yZXYueG1sRE+7CsIw&#xA;FN0F/yFcwU1TCz6oRhGh4iZWF7drc22LzU1pota/N4PgeDjv1aYztXhR6yrLCibjCARxbnXFhYLL&#xA;OR0tQDiPrLG2TAo+5GCz7vdWmGj75hO9Ml+IEMIuQQWl900ipctLMujGtiEO3N22Bn2AbSF1i+8Q&#xA;bmoZR9FMGqw4NJTY0K6k/JE9jYLH9TJN98edPtfZVt+K1F9vd63UcNBtlyA8df4v/rkPWkE8j8Pc&#xA;8CY8Abn+AgAA//8DAFBLAQItABQABgAIAAAAIQDb4fbL7gAAAIUBAAATAAAAAAAAAAAAAAAAAAAA&#xA;AABbQ29udGVudF9UeXBlc10ueG1sUEsBAi0AFAAGAAgAAAAhAFr0LFu/AAAAFQEAAAsAAAAAAAAA&#xA;AAAAAAAAHwEAAF9yZWxzLy5yZWxzUEsBAi0AFAAGAAgAAAAhAI0a0Nq+AAAA3QAAAA8AAAAAAAAA&#xA;AAAAAAAABwIAAGRycy9kb3ducmV2LnhtbFBLBQYAAAAAAwADALcAAADyAgAAAAA=&#xA;" strokeweight="2pt"/><v:line id="Line 940" o:spid="_x0000_s1981" style="position:absolute;visibility:visible" from="1137,14635" to="5116,14636" o:connectortype="straight" o:gfxdata="UEsDBBQABgAIAAAAIQDb4fbL7gAAAIUBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbHyQz07DMAyH&#xA;70i8Q+QralM4IITa7kDhCAiNB7ASt43WOlEcyvb2pNu4IODoPz9/n1xv9vOkForiPDdwXVagiI23&#xA;jocG3rdPxR0oScgWJ8/UwIEENu3lRb09BBKV0ywNjCmFe63FjDSjlD4Q50nv44wpl3HQAc0OB9I3&#xA;VXWrjedEnIq03oC27qjHjympx31un0wiTQLq4bS4shrAECZnMGVTvbD9QSnOhDInjzsyuiBXWQP0&#xA;r4R18jfgnHvJr4nOknrFmJ5xzhraRtHWf3Kkpfz/yGo5S+H73hkquyhdjr3R8m2lj09svwAAAP//&#xA;AwBQSwMEFAAGAAgAAAAhAFr0LFu/AAAAFQEAAAsAAABfcmVscy8ucmVsc2zPwWrDMAwG4Ptg72B0&#xA;X5TuUMaI01uh19I+gLGVxCy2jGSy9e1nemrHjpL4P0nD4SetZiPRyNnCruvBUPYcYp4tXC/Htw8w&#xA;Wl0ObuVMFm6kcBhfX4Yzra62kC6xqGlKVgtLreUTUf1CyWnHhXKbTCzJ1VbKjMX5LzcTvvf9HuXR&#xA;gPHJNKdgQU5hB+ZyK23zHztFL6w81c5zQp6m6P9TMfB3PtPWFCczVQtB9N4U2rp2HOA44NMz4y8A&#xA;AAD//wMAUEsDBBQABgAIAAAAIQDiVnVBxQAAAN0AAAAPAAAAZHJzL2Rvd25yZXYueG1sRI9Ba8JA&#xA;FITvBf/D8oTemo2BthqziggRb6XRS27P7DMJyb4N2VXjv+8WCj0OM/MNk20n04s7ja61rGARxSCI&#xA;K6tbrhWcT/nbEoTzyBp7y6TgSQ62m9lLhqm2D/6me+FrESDsUlTQeD+kUrqqIYMusgNx8K52NOiD&#xA;HGupR3wEuOllEscf0mDLYaHBgfYNVV1xMwq68vyeH772+tQXO32pc19erlqp1/m0W4PwNPn/8F/7&#xA;qBUkn8kKft+EJyA3PwAAAP//AwBQSwECLQAUAAYACAAAACEA2+H2y+4AAACFAQAAEwAAAAAAAAAA&#xA;AAAAAAAAAAAAW0NvbnRlbnRfVHlwZXNdLnhtbFBLAQItABQABgAIAAAAIQBa9CxbvwAAABUBAAAL&#xA;AAAAAAAAAAAAAAAAAB8BAABfcmVscy8ucmVsc1BLAQItABQABgAIAAAAIQDiVnVBxQAAAN0AAAAP&#xA;AAAAAAAAAAAAAAAAAAcCAABkcnMvZG93bnJldi54bWxQSwUGAAAAAAMAAwC3AAAA+QIAAAAA&#xA;" strokeweight="2pt"/><v:line id="Line 941" o:spid="_x0000_s1982" style="position:absolute;visibility:visible" from="1129,14360" to="5108,14361" o:connectortype="straight" o:gfxdata="UEsDBBQABgAIAAAAIQDb4fbL7gAAAIUBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbHyQz07DMAyH&#xA;70i8Q+QralM4IITa7kDhCAiNB7ASt43WOlEcyvb2pNu4IODoPz9/n1xv9vOkForiPDdwXVagiI23&#xA;jocG3rdPxR0oScgWJ8/UwIEENu3lRb09BBKV0ywNjCmFe63FjDSjlD4Q50nv44wpl3HQAc0OB9I3&#xA;VXWrjedEnIq03oC27qjHjympx31un0wiTQLq4bS4shrAECZnMGVTvbD9QSnOhDInjzsyuiBXWQP0&#xA;r4R18jfgnHvJr4nOknrFmJ5xzhraRtHWf3Kkpfz/yGo5S+H73hkquyhdjr3R8m2lj09svwAAAP//&#xA;AwBQSwMEFAAGAAgAAAAhAFr0LFu/AAAAFQEAAAsAAABfcmVscy8ucmVsc2zPwWrDMAwG4Ptg72B0&#xA;X5TuUMaI01uh19I+gLGVxCy2jGSy9e1nemrHjpL4P0nD4SetZiPRyNnCruvBUPYcYp4tXC/Htw8w&#xA;Wl0ObuVMFm6kcBhfX4Yzra62kC6xqGlKVgtLreUTUf1CyWnHhXKbTCzJ1VbKjMX5LzcTvvf9HuXR&#xA;gPHJNKdgQU5hB+ZyK23zHztFL6w81c5zQp6m6P9TMfB3PtPWFCczVQtB9N4U2rp2HOA44NMz4y8A&#xA;AAD//wMAUEsDBBQABgAIAAAAIQB8cCzrwgAAAN0AAAAPAAAAZHJzL2Rvd25yZXYueG1sRE/NagIx&#xA;EL4XfIcwQm81q4LV1SiiFio9lKoPMG7GzepmsiSprj69ORR6/Pj+Z4vW1uJKPlSOFfR7GQjiwumK&#xA;SwWH/cfbGESIyBprx6TgTgEW887LDHPtbvxD110sRQrhkKMCE2OTSxkKQxZDzzXEiTs5bzEm6Eup&#xA;Pd5SuK3lIMtG0mLFqcFgQytDxWX3axVs/fHr0n+URh556zf193oS7Fmp1267nIKI1MZ/8Z/7UysY&#xA;vA/T/vQmPQE5fwIAAP//AwBQSwECLQAUAAYACAAAACEA2+H2y+4AAACFAQAAEwAAAAAAAAAAAAAA&#xA;AAAAAAAAW0NvbnRlbnRfVHlwZXNdLnhtbFBLAQItABQABgAIAAAAIQBa9CxbvwAAABUBAAALAAAA&#xA;AAAAAAAAAAAAAB8BAABfcmVscy8ucmVsc1BLAQItABQABgAIAAAAIQB8cCzrwgAAAN0AAAAPAAAA&#xA;AAAAAAAAAAAAAAcCAABkcnMvZG93bnJldi54bWxQSwUGAAAAAAMAAwC3AAAA9gIAAAAA&#xA;" strokeweight="1pt"/><v:line id="Line 942" o:spid="_x0000_s1983" style="position:absolute;visibility:visible" from="1129,15455" to="5108,15456" o:connectortype="straight" o:gfxdata="UEsDBBQABgAIAAAAIQDb4fbL7gAAAIUBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbHyQz07DMAyH&#xA;70i8Q+QralM4IITa7kDhCAiNB7ASt43WOlEcyvb2pNu4IODoPz9/n1xv9vOkForiPDdwXVagiI23&#xA;jocG3rdPxR0oScgWJ8/UwIEENu3lRb09BBKV0ywNjCmFe63FjDSjlD4Q50nv44wpl3HQAc0OB9I3&#xA;VXWrjedEnIq03oC27qjHjympx31un0wiTQLq4bS4shrAECZnMGVTvbD9QSnOhDInjzsyuiBXWQP0&#xA;r4R18jfgnHvJr4nOknrFmJ5xzhraRtHWf3Kkpfz/yGo5S+H73hkquyhdjr3R8m2lj09svwAAAP//&#xA;AwBQSwMEFAAGAAgAAAAhAFr0LFu/AAAAFQEAAAsAAABfcmVscy8ucmVsc2zPwWrDMAwG4Ptg72B0&#xA;X5TuUMaI01uh19I+gLGVxCy2jGSy9e1nemrHjpL4P0nD4SetZiPRyNnCruvBUPYcYp4tXC/Htw8w&#xA;Wl0ObuVMFm6kcBhfX4Yzra62kC6xqGlKVgtLreUTUf1CyWnHhXKbTCzJ1VbKjMX5LzcTvvf9HuXR&#xA;gPHJNKdgQU5hB+ZyK23zHztFL6w81c5zQp6m6P9TMfB3PtPWFCczVQtB9N4U2rp2HOA44NMz4y8A&#xA;AAD//wMAUEsDBBQABgAIAAAAIQATPIlwxgAAAN0AAAAPAAAAZHJzL2Rvd25yZXYueG1sRI/dagIx&#xA;FITvBd8hHMG7ml0FbVejlP6A4oXU+gDHzelm6+ZkSVLd9umNUPBymJlvmMWqs404kw+1YwX5KANB&#xA;XDpdc6Xg8Pn+8AgiRGSNjWNS8EsBVst+b4GFdhf+oPM+ViJBOBSowMTYFlKG0pDFMHItcfK+nLcY&#xA;k/SV1B4vCW4bOc6yqbRYc1ow2NKLofK0/7EKNv64PeV/lZFH3vi3Zvf6FOy3UsNB9zwHEamL9/B/&#xA;e60VjGeTHG5v0hOQyysAAAD//wMAUEsBAi0AFAAGAAgAAAAhANvh9svuAAAAhQEAABMAAAAAAAAA&#xA;AAAAAAAAAAAAAFtDb250ZW50X1R5cGVzXS54bWxQSwECLQAUAAYACAAAACEAWvQsW78AAAAVAQAA&#xA;CwAAAAAAAAAAAAAAAAAfAQAAX3JlbHMvLnJlbHNQSwECLQAUAAYACAAAACEAEzyJcMYAAADdAAAA&#xA;DwAAAAAAAAAAAAAAAAAHAgAAZHJzL2Rvd25yZXYueG1sUEsFBgAAAAADAAMAtwAAAPoCAAAAAA==&#xA;" strokeweight="1pt"/><v:line id="Line 943" o:spid="_x0000_s1984" style="position:absolute;visibility:visible" from="1129,15180" to="5108,15180" o:connectortype="straight" o:gfxdata="UEsDBBQABgAIAAAAIQDb4fbL7gAAAIUBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbHyQz07DMAyH&#xA;70i8Q+QralM4IITa7kDhCAiNB7ASt43WOlEcyvb2pNu4IODoPz9/n1xv9vOkForiPDdwXVagiI23&#xA;jocG3rdPxR0oScgWJ8/UwIEENu3lRb09BBKV0ywNjCmFe63FjDSjlD4Q50nv44wpl3HQAc0OB9I3&#xA;VXWrjedEnIq03oC27qjHjympx31un0wiTQLq4bS4shrAECZnMGVTvbD9QSnOhDInjzsyuiBXWQP0&#xA;r4R18jfgnHvJr4nOknrFmJ5xzhraRtHWf3Kkpfz/yGo5S+H73hkquyhdjr3R8m2lj09svwAAAP//&#xA;AwBQSwMEFAAGAAgAAAAhAFr0LFu/AAAAFQEAAAsAAABfcmVscy8ucmVsc2zPwWrDMAwG4Ptg72B0&#xA;X5TuUMaI01uh19I+gLGVxCy2jGSy9e1nemrHjpL4P0nD4SetZiPRyNnCruvBUPYcYp4tXC/Htw8w&#xA;Wl0ObuVMFm6kcBhfX4Yzra62kC6xqGlKVgtLreUTUf1CyWnHhXKbTCzJ1VbKjMX5LzcTvvf9HuXR&#xA;gPHJNKdgQU5hB+ZyK23zHztFL6w81c5zQp6m6P9TMfB3PtPWFCczVQtB9N4U2rp2HOA44NMz4y8A&#xA;AAD//wMAUEsDBBQABgAIAAAAIQDj7hcHxgAAAN0AAAAPAAAAZHJzL2Rvd25yZXYueG1sRI/dagIx&#xA;FITvC75DOELvNOsW+rMaRfyBSi+k1gc4bk43WzcnSxJ19embgtDLYWa+YSazzjbiTD7UjhWMhhkI&#xA;4tLpmisF+6/14BVEiMgaG8ek4EoBZtPewwQL7S78SeddrESCcChQgYmxLaQMpSGLYeha4uR9O28x&#xA;JukrqT1eEtw2Ms+yZ2mx5rRgsKWFofK4O1kFG3/4OI5ulZEH3vhVs12+Bfuj1GO/m49BROrif/je&#xA;ftcK8penHP7epCcgp78AAAD//wMAUEsBAi0AFAAGAAgAAAAhANvh9svuAAAAhQEAABMAAAAAAAAA&#xA;AAAAAAAAAAAAAFtDb250ZW50X1R5cGVzXS54bWxQSwECLQAUAAYACAAAACEAWvQsW78AAAAVAQAA&#xA;CwAAAAAAAAAAAAAAAAAfAQAAX3JlbHMvLnJlbHNQSwECLQAUAAYACAAAACEA4+4XB8YAAADdAAAA&#xA;DwAAAAAAAAAAAAAAAAAHAgAAZHJzL2Rvd25yZXYueG1sUEsFBgAAAAADAAMAtwAAAPoCAAAAAA==&#xA;" strokeweight="1pt"/><v:group id="Group 944" o:spid="_x0000_s1985" style="position:absolute;left:1144;top:14935;width:2507;height:240" coordsize="19999,20000" o:gfxdata="UEsDBBQABgAIAAAAIQDb4fbL7gAAAIUBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbHyQz07DMAyH&#xA;70i8Q+QralM4IITa7kDhCAiNB7ASt43WOlEcyvb2pNu4IODoPz9/n1xv9vOkForiPDdwXVagiI23&#xA;jocG3rdPxR0oScgWJ8/UwIEENu3lRb09BBKV0ywNjCmFe63FjDSjlD4Q50nv44wpl3HQAc0OB9I3&#xA;VXWrjedEnIq03oC27qjHjympx31un0wiTQLq4bS4shrAECZnMGVTvbD9QSnOhDInjzsyuiBXWQP0&#xA;r4R18jfgnHvJr4nOknrFmJ5xzhraRtHWf3Kkpfz/yGo5S+H73hkquyhdjr3R8m2lj09svwAAAP//&#xA;AwBQSwMEFAAGAAgAAAAhAFr0LFu/AAAAFQEAAAsAAABfcmVscy8ucmVsc2zPwWrDMAwG4Ptg72B0&#xA;X5TuUMaI01uh19I+gLGVxCy2jGSy9e1nemrHjpL4P0nD4SetZiPRyNnCruvBUPYcYp4tXC/Htw8w&#xA;Wl0ObuVMFm6kcBhfX4Yzra62kC6xqGlKVgtLreUTUf1CyWnHhXKbTCzJ1VbKjMX5LzcTvvf9HuXR&#xA;gPHJNKdgQU5hB+ZyK23zHztFL6w81c5zQp6m6P9TMfB3PtPWFCczVQtB9N4U2rp2HOA44NMz4y8A&#xA;AAD//wMAUEsDBBQABgAIAAAAIQDo5I3dxgAAAN0AAAAPAAAAZHJzL2Rvd25yZXYueG1sRI9Ba8JA&#xA;FITvBf/D8oTe6iaGthJdRURLDyJUBfH2yD6TYPZtyK5J/PeuIPQ4zMw3zGzRm0q01LjSsoJ4FIEg&#xA;zqwuOVdwPGw+JiCcR9ZYWSYFd3KwmA/eZphq2/EftXufiwBhl6KCwvs6ldJlBRl0I1sTB+9iG4M+&#xA;yCaXusEuwE0lx1H0JQ2WHBYKrGlVUHbd34yCnw67ZRKv2+31srqfD5+70zYmpd6H/XIKwlPv/8Ov&#xA;9q9WMP5OEni+CU9Azh8AAAD//wMAUEsBAi0AFAAGAAgAAAAhANvh9svuAAAAhQEAABMAAAAAAAAA&#xA;AAAAAAAAAAAAAFtDb250ZW50X1R5cGVzXS54bWxQSwECLQAUAAYACAAAACEAWvQsW78AAAAVAQAA&#xA;CwAAAAAAAAAAAAAAAAAfAQAAX3JlbHMvLnJlbHNQSwECLQAUAAYACAAAACEA6OSN3cYAAADdAAAA&#xA;DwAAAAAAAAAAAAAAAAAHAgAAZHJzL2Rvd25yZXYueG1sUEsFBgAAAAADAAMAtwAAAPoCAAAAAA==&#xA;"><v:rect id="Rectangle 945" o:spid="_x0000_s1986" style="position:absolute;width:8856;height:20000;visibility:visible" o:gfxdata="UEsDBBQABgAIAAAAIQDb4fbL7gAAAIUBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbHyQz07DMAyH&#xA;70i8Q+QralM4IITa7kDhCAiNB7ASt43WOlEcyvb2pNu4IODoPz9/n1xv9vOkForiPDdwXVagiI23&#xA;jocG3rdPxR0oScgWJ8/UwIEENu3lRb09BBKV0ywNjCmFe63FjDSjlD4Q50nv44wpl3HQAc0OB9I3&#xA;VXWrjedEnIq03oC27qjHjympx31un0wiTQLq4bS4shrAECZnMGVTvbD9QSnOhDInjzsyuiBXWQP0&#xA;r4R18jfgnHvJr4nOknrFmJ5xzhraRtHWf3Kkpfz/yGo5S+H73hkquyhdjr3R8m2lj09svwAAAP//&#xA;AwBQSwMEFAAGAAgAAAAhAFr0LFu/AAAAFQEAAAsAAABfcmVscy8ucmVsc2zPwWrDMAwG4Ptg72B0&#xA;X5TuUMaI01uh19I+gLGVxCy2jGSy9e1nemrHjpL4P0nD4SetZiPRyNnCruvBUPYcYp4tXC/Htw8w&#xA;Wl0ObuVMFm6kcBhfX4Yzra62kC6xqGlKVgtLreUTUf1CyWnHhXKbTCzJ1VbKjMX5LzcTvvf9HuXR&#xA;gPHJNKdgQU5hB+ZyK23zHztFL6w81c5zQp6m6P9TMfB3PtPWFCczVQtB9N4U2rp2HOA44NMz4y8A&#xA;AAD//wMAUEsDBBQABgAIAAAAIQAuwq90wgAAAN0AAAAPAAAAZHJzL2Rvd25yZXYueG1sRI9Bi8Iw&#xA;FITvgv8hvAVvmq6K61ajFEHwanVhj4/m2Xa3ealJ1PrvjSB4HGbmG2a57kwjruR8bVnB5ygBQVxY&#xA;XXOp4HjYDucgfEDW2FgmBXfysF71e0tMtb3xnq55KEWEsE9RQRVCm0rpi4oM+pFtiaN3ss5giNKV&#xA;Uju8Rbhp5DhJZtJgzXGhwpY2FRX/+cUoyLK/7uecf+PWy3niZnqqy+xXqcFHly1ABOrCO/xq77SC&#xA;8ddkCs838QnI1QMAAP//AwBQSwECLQAUAAYACAAAACEA2+H2y+4AAACFAQAAEwAAAAAAAAAAAAAA&#xA;AAAAAAAAW0NvbnRlbnRfVHlwZXNdLnhtbFBLAQItABQABgAIAAAAIQBa9CxbvwAAABUBAAALAAAA&#xA;AAAAAAAAAAAAAB8BAABfcmVscy8ucmVsc1BLAQItABQABgAIAAAAIQAuwq90wgAAAN0AAAAPAAAA&#xA;AAAAAAAAAAAAAAcCAABkcnMvZG93bnJldi54bWxQSwUGAAAAAAMAAwC3AAAA9gIAAAAA&#xA;" filled="f" stroked="f" strokeweight=".25pt"><v:textbox style="mso-next-textbox:#Rectangle 945" inset="1pt,1pt,1pt,1pt"><w:txbxContent><w:p w14:paraId="5800332B" w14:textId="77777777" w:rsidR="00B827FD" w:rsidRPr="00951A0B" w:rsidRDefault="00B827FD" w:rsidP="00B827FD"><w:pPr><w:pStyle w:val="a1"/><w:jc w:val="left"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr><w:r w:rsidRPr="00951A0B"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>Разраб.</w:t></w:r></w:p></w:txbxContent></v:textbox></v:rect><v:rect id="Rectangle 946" o:spid="_x0000_s1987" style="position:absolute;left:9281;width:10718;height:20000;visibility:visible" o:gfxdata="UEsDBBQABgAIAAAAIQDb4fbL7gAAAIUBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbHyQz07DMAyH&#xA;70i8Q+QralM4IITa7kDhCAiNB7ASt43WOlEcyvb2pNu4IODoPz9/n1xv9vOkForiPDdwXVagiI23&#xA;jocG3rdPxR0oScgWJ8/UwIEENu3lRb09BBKV0ywNjCmFe63FjDSjlD4Q50nv44wpl3HQAc0OB9I3&#xA;VXWrjedEnIq03oC27qjHjympx31un0wiTQLq4bS4shrAECZnMGVTvbD9QSnOhDInjzsyuiBXWQP0&#xA;r4R18jfgnHvJr4nOknrFmJ5xzhraRtHWf3Kkpfz/yGo5S+H73hkquyhdjr3R8m2lj09svwAAAP//&#xA;AwBQSwMEFAAGAAgAAAAhAFr0LFu/AAAAFQEAAAsAAABfcmVscy8ucmVsc2zPwWrDMAwG4Ptg72B0&#xA;X5TuUMaI01uh19I+gLGVxCy2jGSy9e1nemrHjpL4P0nD4SetZiPRyNnCruvBUPYcYp4tXC/Htw8w&#xA;Wl0ObuVMFm6kcBhfX4Yzra62kC6xqGlKVgtLreUTUf1CyWnHhXKbTCzJ1VbKjMX5LzcTvvf9HuXR&#xA;gPHJNKdgQU5hB+ZyK23zHztFL6w81c5zQp6m6P9TMfB3PtPWFCczVQtB9N4U2rp2HOA44NMz4y8A&#xA;AAD//wMAUEsDBBQABgAIAAAAIQBBjgrvxAAAAN0AAAAPAAAAZHJzL2Rvd25yZXYueG1sRI9Ba8JA&#xA;FITvgv9heUJvdaNtU02zkVAQem1U6PGRfSap2bdxd9X033cLBY/DzHzD5JvR9OJKzneWFSzmCQji&#xA;2uqOGwX73fZxBcIHZI29ZVLwQx42xXSSY6btjT/pWoVGRAj7DBW0IQyZlL5uyaCf24E4ekfrDIYo&#xA;XSO1w1uEm14ukySVBjuOCy0O9N5SfaouRkFZfo+Hc7XGrZerxKX6WTfll1IPs7F8AxFoDPfwf/tD&#xA;K1i+Pr3A35v4BGTxCwAA//8DAFBLAQItABQABgAIAAAAIQDb4fbL7gAAAIUBAAATAAAAAAAAAAAA&#xA;AAAAAAAAAABbQ29udGVudF9UeXBlc10ueG1sUEsBAi0AFAAGAAgAAAAhAFr0LFu/AAAAFQEAAAsA&#xA;AAAAAAAAAAAAAAAAHwEAAF9yZWxzLy5yZWxzUEsBAi0AFAAGAAgAAAAhAEGOCu/EAAAA3QAAAA8A&#xA;AAAAAAAAAAAAAAAABwIAAGRycy9kb3ducmV2LnhtbFBLBQYAAAAAAwADALcAAAD4AgAAAAA=&#xA;" filled="f" stroked="f" strokeweight=".25pt"><v:textbox style="mso-next-textbox:#Rectangle 946" inset="1pt,1pt,1pt,1pt"><w:txbxContent><w:p w14:paraId="1AB7E379" w14:textId="461DE9EC" w:rsidR="00B827FD" w:rsidRPr="00017D2F" w:rsidRDefault="00B827FD" w:rsidP="00B827FD"><w:pPr><w:pStyle w:val="a1"/><w:jc w:val="left"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="ru-RU"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="ru-RU"/></w:rPr><w:t>Мацуев И.М.</w:t></w:r></w:p></w:txbxContent></v:textbox></v:rect></v:group><v:group id="Group 947" o:spid="_x0000_s1988" style="position:absolute;left:1144;top:15204;width:2507;height:239" coordsize="19999,20000" o:gfxdata="UEsDBBQABgAIAAAAIQDb4fbL7gAAAIUBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbHyQz07DMAyH&#xA;70i8Q+QralM4IITa7kDhCAiNB7ASt43WOlEcyvb2pNu4IODoPz9/n1xv9vOkForiPDdwXVagiI23&#xA;jocG3rdPxR0oScgWJ8/UwIEENu3lRb09BBKV0ywNjCmFe63FjDSjlD4Q50nv44wpl3HQAc0OB9I3&#xA;VXWrjedEnIq03oC27qjHjympx31un0wiTQLq4bS4shrAECZnMGVTvbD9QSnOhDInjzsyuiBXWQP0&#xA;r4R18jfgnHvJr4nOknrFmJ5xzhraRtHWf3Kkpfz/yGo5S+H73hkquyhdjr3R8m2lj09svwAAAP//&#xA;AwBQSwMEFAAGAAgAAAAhAFr0LFu/AAAAFQEAAAsAAABfcmVscy8ucmVsc2zPwWrDMAwG4Ptg72B0&#xA;X5TuUMaI01uh19I+gLGVxCy2jGSy9e1nemrHjpL4P0nD4SetZiPRyNnCruvBUPYcYp4tXC/Htw8w&#xA;Wl0ObuVMFm6kcBhfX4Yzra62kC6xqGlKVgtLreUTUf1CyWnHhXKbTCzJ1VbKjMX5LzcTvvf9HuXR&#xA;gPHJNKdgQU5hB+ZyK23zHztFL6w81c5zQp6m6P9TMfB3PtPWFCczVQtB9N4U2rp2HOA44NMz4y8A&#xA;AAD//wMAUEsDBBQABgAIAAAAIQD4ky5FxgAAAN0AAAAPAAAAZHJzL2Rvd25yZXYueG1sRI9Pi8Iw&#xA;FMTvC/sdwlvwpmkVdalGEdFlDyL4BxZvj+bZFpuX0sS2fnsjCHscZuY3zHzZmVI0VLvCsoJ4EIEg&#xA;Tq0uOFNwPm373yCcR9ZYWiYFD3KwXHx+zDHRtuUDNUefiQBhl6CC3PsqkdKlORl0A1sRB+9qa4M+&#xA;yDqTusY2wE0ph1E0kQYLDgs5VrTOKb0d70bBT4vtahRvmt3tun5cTuP93y4mpXpf3WoGwlPn/8Pv&#xA;9q9WMJyOJvB6E56AXDwBAAD//wMAUEsBAi0AFAAGAAgAAAAhANvh9svuAAAAhQEAABMAAAAAAAAA&#xA;AAAAAAAAAAAAAFtDb250ZW50X1R5cGVzXS54bWxQSwECLQAUAAYACAAAACEAWvQsW78AAAAVAQAA&#xA;CwAAAAAAAAAAAAAAAAAfAQAAX3JlbHMvLnJlbHNQSwECLQAUAAYACAAAACEA+JMuRcYAAADdAAAA&#xA;DwAAAAAAAAAAAAAAAAAHAgAAZHJzL2Rvd25yZXYueG1sUEsFBgAAAAADAAMAtwAAAPoCAAAAAA==&#xA;"><v:rect id="Rectangle 948" o:spid="_x0000_s1989" style="position:absolute;width:8856;height:20000;visibility:visible" o:gfxdata="UEsDBBQABgAIAAAAIQDb4fbL7gAAAIUBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbHyQz07DMAyH&#xA;70i8Q+QralM4IITa7kDhCAiNB7ASt43WOlEcyvb2pNu4IODoPz9/n1xv9vOkForiPDdwXVagiI23&#xA;jocG3rdPxR0oScgWJ8/UwIEENu3lRb09BBKV0ywNjCmFe63FjDSjlD4Q50nv44wpl3HQAc0OB9I3&#xA;VXWrjedEnIq03oC27qjHjympx31un0wiTQLq4bS4shrAECZnMGVTvbD9QSnOhDInjzsyuiBXWQP0&#xA;r4R18jfgnHvJr4nOknrFmJ5xzhraRtHWf3Kkpfz/yGo5S+H73hkquyhdjr3R8m2lj09svwAAAP//&#xA;AwBQSwMEFAAGAAgAAAAhAFr0LFu/AAAAFQEAAAsAAABfcmVscy8ucmVsc2zPwWrDMAwG4Ptg72B0&#xA;X5TuUMaI01uh19I+gLGVxCy2jGSy9e1nemrHjpL4P0nD4SetZiPRyNnCruvBUPYcYp4tXC/Htw8w&#xA;Wl0ObuVMFm6kcBhfX4Yzra62kC6xqGlKVgtLreUTUf1CyWnHhXKbTCzJ1VbKjMX5LzcTvvf9HuXR&#xA;gPHJNKdgQU5hB+ZyK23zHztFL6w81c5zQp6m6P9TMfB3PtPWFCczVQtB9N4U2rp2HOA44NMz4y8A&#xA;AAD//wMAUEsDBBQABgAIAAAAIQDeEDEDwwAAAN0AAAAPAAAAZHJzL2Rvd25yZXYueG1sRI9Pi8Iw&#xA;FMTvC36H8ARva6ou/qlGKYLgdbsKHh/Ns602LzWJWr+9WVjY4zAzv2FWm8404kHO15YVjIYJCOLC&#xA;6ppLBYef3ecchA/IGhvLpOBFHjbr3scKU22f/E2PPJQiQtinqKAKoU2l9EVFBv3QtsTRO1tnMETp&#xA;SqkdPiPcNHKcJFNpsOa4UGFL24qKa343CrLs0h1v+QJ3Xs4TN9VfusxOSg36XbYEEagL/+G/9l4r&#xA;GM8mM/h9E5+AXL8BAAD//wMAUEsBAi0AFAAGAAgAAAAhANvh9svuAAAAhQEAABMAAAAAAAAAAAAA&#xA;AAAAAAAAAFtDb250ZW50X1R5cGVzXS54bWxQSwECLQAUAAYACAAAACEAWvQsW78AAAAVAQAACwAA&#xA;AAAAAAAAAAAAAAAfAQAAX3JlbHMvLnJlbHNQSwECLQAUAAYACAAAACEA3hAxA8MAAADdAAAADwAA&#xA;AAAAAAAAAAAAAAAHAgAAZHJzL2Rvd25yZXYueG1sUEsFBgAAAAADAAMAtwAAAPcCAAAAAA==&#xA;" filled="f" stroked="f" strokeweight=".25pt"><v:textbox style="mso-next-textbox:#Rectangle 948" inset="1pt,1pt,1pt,1pt"><w:txbxContent><w:p w14:paraId="1B9E4D2F" w14:textId="77777777" w:rsidR="00B827FD" w:rsidRPr="00951A0B" w:rsidRDefault="00B827FD" w:rsidP="00B827FD"><w:pPr><w:pStyle w:val="a1"/><w:jc w:val="left"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr><w:r w:rsidRPr="00951A0B"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>Провер.</w:t></w:r></w:p></w:txbxContent></v:textbox></v:rect><v:rect id="Rectangle 949" o:spid="_x0000_s1990" style="position:absolute;left:9281;width:10718;height:20000;visibility:visible" o:gfxdata="UEsDBBQABgAIAAAAIQDb4fbL7gAAAIUBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbHyQz07DMAyH&#xA;70i8Q+QralM4IITa7kDhCAiNB7ASt43WOlEcyvb2pNu4IODoPz9/n1xv9vOkForiPDdwXVagiI23&#xA;jocG3rdPxR0oScgWJ8/UwIEENu3lRb09BBKV0ywNjCmFe63FjDSjlD4Q50nv44wpl3HQAc0OB9I3&#xA;VXWrjedEnIq03oC27qjHjympx31un0wiTQLq4bS4shrAECZnMGVTvbD9QSnOhDInjzsyuiBXWQP0&#xA;r4R18jfgnHvJr4nOknrFmJ5xzhraRtHWf3Kkpfz/yGo5S+H73hkquyhdjr3R8m2lj09svwAAAP//&#xA;AwBQSwMEFAAGAAgAAAAhAFr0LFu/AAAAFQEAAAsAAABfcmVscy8ucmVsc2zPwWrDMAwG4Ptg72B0&#xA;X5TuUMaI01uh19I+gLGVxCy2jGSy9e1nemrHjpL4P0nD4SetZiPRyNnCruvBUPYcYp4tXC/Htw8w&#xA;Wl0ObuVMFm6kcBhfX4Yzra62kC6xqGlKVgtLreUTUf1CyWnHhXKbTCzJ1VbKjMX5LzcTvvf9HuXR&#xA;gPHJNKdgQU5hB+ZyK23zHztFL6w81c5zQp6m6P9TMfB3PtPWFCczVQtB9N4U2rp2HOA44NMz4y8A&#xA;AAD//wMAUEsDBBQABgAIAAAAIQCvj6VxwAAAAN0AAAAPAAAAZHJzL2Rvd25yZXYueG1sRE/LisIw&#xA;FN0P+A/hCu7G1Afa6RilDAhurQqzvDR32mpzU5OM1r83C8Hl4bxXm9604kbON5YVTMYJCOLS6oYr&#xA;BcfD9jMF4QOyxtYyKXiQh8168LHCTNs77+lWhErEEPYZKqhD6DIpfVmTQT+2HXHk/qwzGCJ0ldQO&#xA;7zHctHKaJAtpsOHYUGNHPzWVl+LfKMjzc3+6Fl+49TJN3ELPdZX/KjUa9vk3iEB9eItf7p1WMF3O&#xA;4tz4Jj4BuX4CAAD//wMAUEsBAi0AFAAGAAgAAAAhANvh9svuAAAAhQEAABMAAAAAAAAAAAAAAAAA&#xA;AAAAAFtDb250ZW50X1R5cGVzXS54bWxQSwECLQAUAAYACAAAACEAWvQsW78AAAAVAQAACwAAAAAA&#xA;AAAAAAAAAAAfAQAAX3JlbHMvLnJlbHNQSwECLQAUAAYACAAAACEAr4+lccAAAADdAAAADwAAAAAA&#xA;AAAAAAAAAAAHAgAAZHJzL2Rvd25yZXYueG1sUEsFBgAAAAADAAMAtwAAAPQCAAAAAA==&#xA;" filled="f" stroked="f" strokeweight=".25pt"><v:textbox style="mso-next-textbox:#Rectangle 949" inset="1pt,1pt,1pt,1pt"><w:txbxContent><w:p w14:paraId="2303CEB4" w14:textId="77777777" w:rsidR="00B827FD" w:rsidRPr="00CF1877" w:rsidRDefault="00B827FD" w:rsidP="00B827FD"><w:pPr><w:rPr><w:i/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="ru-RU"/></w:rPr></w:pPr><w:r><w:rPr><w:i/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="ru-RU"/></w:rPr><w:t>Северинчик Н.А.</w:t></w:r></w:p><w:p w14:paraId="641DA5EE" w14:textId="77777777" w:rsidR="00B827FD" w:rsidRPr="00CF1877" w:rsidRDefault="00B827FD" w:rsidP="00B827FD"><w:pPr><w:rPr><w:i/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="ru-RU"/></w:rPr></w:pPr><w:r><w:rPr><w:i/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="ru-RU"/></w:rPr><w:t>.</w:t></w:r></w:p><w:p w14:paraId="607393CB" w14:textId="77777777" w:rsidR="00B827FD" w:rsidRPr="00CB455E" w:rsidRDefault="00B827FD" w:rsidP="00B827FD"><w:pPr><w:rPr><w:szCs w:val="16"/></w:rPr></w:pPr></w:p></w:txbxContent></v:textbox></v:rect></v:group><v:group id="Group 950" o:spid="_x0000_s1991" style="position:absolute;left:1144;top:15479;width:2507;height:239" coordsize="19999,20000" o:gfxdata="UEsDBBQABgAIAAAAIQDb4fbL7gAAAIUBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbHyQz07DMAyH&#xA;70i8Q+QralM4IITa7kDhCAiNB7ASt43WOlEcyvb2pNu4IODoPz9/n1xv9vOkForiPDdwXVagiI23&#xA;jocG3rdPxR0oScgWJ8/UwIEENu3lRb09BBKV0ywNjCmFe63FjDSjlD4Q50nv44wpl3HQAc0OB9I3&#xA;VXWrjedEnIq03oC27qjHjympx31un0wiTQLq4bS4shrAECZnMGVTvbD9QSnOhDInjzsyuiBXWQP0&#xA;r4R18jfgnHvJr4nOknrFmJ5xzhraRtHWf3Kkpfz/yGo5S+H73hkquyhdjr3R8m2lj09svwAAAP//&#xA;AwBQSwMEFAAGAAgAAAAhAFr0LFu/AAAAFQEAAAsAAABfcmVscy8ucmVsc2zPwWrDMAwG4Ptg72B0&#xA;X5TuUMaI01uh19I+gLGVxCy2jGSy9e1nemrHjpL4P0nD4SetZiPRyNnCruvBUPYcYp4tXC/Htw8w&#xA;Wl0ObuVMFm6kcBhfX4Yzra62kC6xqGlKVgtLreUTUf1CyWnHhXKbTCzJ1VbKjMX5LzcTvvf9HuXR&#xA;gPHJNKdgQU5hB+ZyK23zHztFL6w81c5zQp6m6P9TMfB3PtPWFCczVQtB9N4U2rp2HOA44NMz4y8A&#xA;AAD//wMAUEsDBBQABgAIAAAAIQCJDLo3xgAAAN0AAAAPAAAAZHJzL2Rvd25yZXYueG1sRI9Ba8JA&#xA;FITvQv/D8gre6iaKtkZXEbHSgwiNBfH2yD6TYPZtyG6T+O+7QsHjMDPfMMt1byrRUuNKywriUQSC&#xA;OLO65FzBz+nz7QOE88gaK8uk4E4O1quXwRITbTv+pjb1uQgQdgkqKLyvEyldVpBBN7I1cfCutjHo&#xA;g2xyqRvsAtxUchxFM2mw5LBQYE3bgrJb+msU7DvsNpN41x5u1+39cpoez4eYlBq+9psFCE+9f4b/&#xA;219awfh9MofHm/AE5OoPAAD//wMAUEsBAi0AFAAGAAgAAAAhANvh9svuAAAAhQEAABMAAAAAAAAA&#xA;AAAAAAAAAAAAAFtDb250ZW50X1R5cGVzXS54bWxQSwECLQAUAAYACAAAACEAWvQsW78AAAAVAQAA&#xA;CwAAAAAAAAAAAAAAAAAfAQAAX3JlbHMvLnJlbHNQSwECLQAUAAYACAAAACEAiQy6N8YAAADdAAAA&#xA;DwAAAAAAAAAAAAAAAAAHAgAAZHJzL2Rvd25yZXYueG1sUEsFBgAAAAADAAMAtwAAAPoCAAAAAA==&#xA;"><v:rect id="Rectangle 951" o:spid="_x0000_s1992" style="position:absolute;width:8856;height:20000;visibility:visible" o:gfxdata="UEsDBBQABgAIAAAAIQDb4fbL7gAAAIUBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbHyQz07DMAyH&#xA;70i8Q+QralM4IITa7kDhCAiNB7ASt43WOlEcyvb2pNu4IODoPz9/n1xv9vOkForiPDdwXVagiI23&#xA;jocG3rdPxR0oScgWJ8/UwIEENu3lRb09BBKV0ywNjCmFe63FjDSjlD4Q50nv44wpl3HQAc0OB9I3&#xA;VXWrjedEnIq03oC27qjHjympx31un0wiTQLq4bS4shrAECZnMGVTvbD9QSnOhDInjzsyuiBXWQP0&#xA;r4R18jfgnHvJr4nOknrFmJ5xzhraRtHWf3Kkpfz/yGo5S+H73hkquyhdjr3R8m2lj09svwAAAP//&#xA;AwBQSwMEFAAGAAgAAAAhAFr0LFu/AAAAFQEAAAsAAABfcmVscy8ucmVsc2zPwWrDMAwG4Ptg72B0&#xA;X5TuUMaI01uh19I+gLGVxCy2jGSy9e1nemrHjpL4P0nD4SetZiPRyNnCruvBUPYcYp4tXC/Htw8w&#xA;Wl0ObuVMFm6kcBhfX4Yzra62kC6xqGlKVgtLreUTUf1CyWnHhXKbTCzJ1VbKjMX5LzcTvvf9HuXR&#xA;gPHJNKdgQU5hB+ZyK23zHztFL6w81c5zQp6m6P9TMfB3PtPWFCczVQtB9N4U2rp2HOA44NMz4y8A&#xA;AAD//wMAUEsDBBQABgAIAAAAIQAJ/9oKwQAAAN0AAAAPAAAAZHJzL2Rvd25yZXYueG1sRE/Pa4Mw&#xA;FL4P+j+EV9itxpbSOmcsUhB2nVthx4d5UzfzYpPUuv9+OQx2/Ph+F6fFjGIm5wfLCrZJCoK4tXrg&#xA;TsH7W73JQPiArHG0TAp+yMOpXD0UmGt751eam9CJGMI+RwV9CFMupW97MugTOxFH7tM6gyFC10nt&#xA;8B7DzSh3aXqQBgeODT1OdO6p/W5uRkFVfS2Xa/OEtZdZ6g56r7vqQ6nH9VI9gwi0hH/xn/tFK9gd&#xA;93F/fBOfgCx/AQAA//8DAFBLAQItABQABgAIAAAAIQDb4fbL7gAAAIUBAAATAAAAAAAAAAAAAAAA&#xA;AAAAAABbQ29udGVudF9UeXBlc10ueG1sUEsBAi0AFAAGAAgAAAAhAFr0LFu/AAAAFQEAAAsAAAAA&#xA;AAAAAAAAAAAAHwEAAF9yZWxzLy5yZWxzUEsBAi0AFAAGAAgAAAAhAAn/2grBAAAA3QAAAA8AAAAA&#xA;AAAAAAAAAAAABwIAAGRycy9kb3ducmV2LnhtbFBLBQYAAAAAAwADALcAAAD1AgAAAAA=&#xA;" filled="f" stroked="f" strokeweight=".25pt"><v:textbox style="mso-next-textbox:#Rectangle 951" inset="1pt,1pt,1pt,1pt"><w:txbxContent><w:p w14:paraId="4CF9C802" w14:textId="77777777" w:rsidR="00B827FD" w:rsidRPr="0042148F" w:rsidRDefault="00B827FD" w:rsidP="00B827FD"><w:pPr><w:pStyle w:val="a1"/><w:jc w:val="left"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="16"/><w:szCs w:val="16"/><w:lang w:val="ru-RU"/></w:rPr></w:pPr><w:r w:rsidRPr="00951A0B"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="ru-RU"/></w:rPr><w:t>Консульт</w:t></w:r><w:r w:rsidRPr="0042148F"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="16"/><w:szCs w:val="16"/><w:lang w:val="ru-RU"/></w:rPr><w:t>.</w:t></w:r></w:p></w:txbxContent></v:textbox></v:rect><v:rect id="Rectangle 952" o:spid="_x0000_s1993" style="position:absolute;left:9281;width:10718;height:20000;visibility:visible" o:gfxdata="UEsDBBQABgAIAAAAIQDb4fbL7gAAAIUBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbHyQz07DMAyH&#xA;70i8Q+QralM4IITa7kDhCAiNB7ASt43WOlEcyvb2pNu4IODoPz9/n1xv9vOkForiPDdwXVagiI23&#xA;jocG3rdPxR0oScgWJ8/UwIEENu3lRb09BBKV0ywNjCmFe63FjDSjlD4Q50nv44wpl3HQAc0OB9I3&#xA;VXWrjedEnIq03oC27qjHjympx31un0wiTQLq4bS4shrAECZnMGVTvbD9QSnOhDInjzsyuiBXWQP0&#xA;r4R18jfgnHvJr4nOknrFmJ5xzhraRtHWf3Kkpfz/yGo5S+H73hkquyhdjr3R8m2lj09svwAAAP//&#xA;AwBQSwMEFAAGAAgAAAAhAFr0LFu/AAAAFQEAAAsAAABfcmVscy8ucmVsc2zPwWrDMAwG4Ptg72B0&#xA;X5TuUMaI01uh19I+gLGVxCy2jGSy9e1nemrHjpL4P0nD4SetZiPRyNnCruvBUPYcYp4tXC/Htw8w&#xA;Wl0ObuVMFm6kcBhfX4Yzra62kC6xqGlKVgtLreUTUf1CyWnHhXKbTCzJ1VbKjMX5LzcTvvf9HuXR&#xA;gPHJNKdgQU5hB+ZyK23zHztFL6w81c5zQp6m6P9TMfB3PtPWFCczVQtB9N4U2rp2HOA44NMz4y8A&#xA;AAD//wMAUEsDBBQABgAIAAAAIQBms3+RwgAAAN0AAAAPAAAAZHJzL2Rvd25yZXYueG1sRI9Bi8Iw&#xA;FITvgv8hPGFvmiqitRqlLAh7tavg8dE822rz0k2y2v33RhD2OMzMN8xm15tW3Mn5xrKC6SQBQVxa&#xA;3XCl4Pi9H6cgfEDW2FomBX/kYbcdDjaYafvgA92LUIkIYZ+hgjqELpPSlzUZ9BPbEUfvYp3BEKWr&#xA;pHb4iHDTylmSLKTBhuNCjR191lTeil+jIM+v/emnWOHeyzRxCz3XVX5W6mPU52sQgfrwH363v7SC&#xA;2XI+hdeb+ATk9gkAAP//AwBQSwECLQAUAAYACAAAACEA2+H2y+4AAACFAQAAEwAAAAAAAAAAAAAA&#xA;AAAAAAAAW0NvbnRlbnRfVHlwZXNdLnhtbFBLAQItABQABgAIAAAAIQBa9CxbvwAAABUBAAALAAAA&#xA;AAAAAAAAAAAAAB8BAABfcmVscy8ucmVsc1BLAQItABQABgAIAAAAIQBms3+RwgAAAN0AAAAPAAAA&#xA;AAAAAAAAAAAAAAcCAABkcnMvZG93bnJldi54bWxQSwUGAAAAAAMAAwC3AAAA9gIAAAAA&#xA;" filled="f" stroked="f" strokeweight=".25pt"><v:textbox style="mso-next-textbox:#Rectangle 952" inset="1pt,1pt,1pt,1pt"><w:txbxContent><w:p w14:paraId="257AD8D5" w14:textId="77777777" w:rsidR="00B827FD" w:rsidRPr="00CF1877" w:rsidRDefault="00B827FD" w:rsidP="00B827FD"><w:pPr><w:rPr><w:i/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="ru-RU"/></w:rPr></w:pPr><w:r><w:rPr><w:i/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="ru-RU"/></w:rPr><w:t>Северинчик Н.А.</w:t></w:r></w:p><w:p w14:paraId="74DD30A1" w14:textId="77777777" w:rsidR="00B827FD" w:rsidRPr="00CB455E" w:rsidRDefault="00B827FD" w:rsidP="00B827FD"><w:pPr><w:rPr><w:szCs w:val="16"/></w:rPr></w:pPr></w:p></w:txbxContent></v:textbox></v:rect></v:group><v:group id="Group 953" o:spid="_x0000_s1994" style="position:absolute;left:1144;top:15746;width:2507;height:241" coordsize="19999,20000" o:gfxdata="UEsDBBQABgAIAAAAIQDb4fbL7gAAAIUBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbHyQz07DMAyH&#xA;70i8Q+QralM4IITa7kDhCAiNB7ASt43WOlEcyvb2pNu4IODoPz9/n1xv9vOkForiPDdwXVagiI23&#xA;jocG3rdPxR0oScgWJ8/UwIEENu3lRb09BBKV0ywNjCmFe63FjDSjlD4Q50nv44wpl3HQAc0OB9I3&#xA;VXWrjedEnIq03oC27qjHjympx31un0wiTQLq4bS4shrAECZnMGVTvbD9QSnOhDInjzsyuiBXWQP0&#xA;r4R18jfgnHvJr4nOknrFmJ5xzhraRtHWf3Kkpfz/yGo5S+H73hkquyhdjr3R8m2lj09svwAAAP//&#xA;AwBQSwMEFAAGAAgAAAAhAFr0LFu/AAAAFQEAAAsAAABfcmVscy8ucmVsc2zPwWrDMAwG4Ptg72B0&#xA;X5TuUMaI01uh19I+gLGVxCy2jGSy9e1nemrHjpL4P0nD4SetZiPRyNnCruvBUPYcYp4tXC/Htw8w&#xA;Wl0ObuVMFm6kcBhfX4Yzra62kC6xqGlKVgtLreUTUf1CyWnHhXKbTCzJ1VbKjMX5LzcTvvf9HuXR&#xA;gPHJNKdgQU5hB+ZyK23zHztFL6w81c5zQp6m6P9TMfB3PtPWFCczVQtB9N4U2rp2HOA44NMz4y8A&#xA;AAD//wMAUEsDBBQABgAIAAAAIQDfrls7xgAAAN0AAAAPAAAAZHJzL2Rvd25yZXYueG1sRI9Ba8JA&#xA;FITvhf6H5RW86SaxtpK6ikgtHkRQC8XbI/tMgtm3Ibsm8d+7gtDjMDPfMLNFbyrRUuNKywriUQSC&#xA;OLO65FzB73E9nIJwHlljZZkU3MjBYv76MsNU24731B58LgKEXYoKCu/rVEqXFWTQjWxNHLyzbQz6&#xA;IJtc6ga7ADeVTKLoQxosOSwUWNOqoOxyuBoFPx12y3H83W4v59XtdJzs/rYxKTV465dfIDz1/j/8&#xA;bG+0guTzPYHHm/AE5PwOAAD//wMAUEsBAi0AFAAGAAgAAAAhANvh9svuAAAAhQEAABMAAAAAAAAA&#xA;AAAAAAAAAAAAAFtDb250ZW50X1R5cGVzXS54bWxQSwECLQAUAAYACAAAACEAWvQsW78AAAAVAQAA&#xA;CwAAAAAAAAAAAAAAAAAfAQAAX3JlbHMvLnJlbHNQSwECLQAUAAYACAAAACEA365bO8YAAADdAAAA&#xA;DwAAAAAAAAAAAAAAAAAHAgAAZHJzL2Rvd25yZXYueG1sUEsFBgAAAAADAAMAtwAAAPoCAAAAAA==&#xA;"><v:rect id="Rectangle 954" o:spid="_x0000_s1995" style="position:absolute;width:8856;height:20000;visibility:visible" o:gfxdata="UEsDBBQABgAIAAAAIQDb4fbL7gAAAIUBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbHyQz07DMAyH&#xA;70i8Q+QralM4IITa7kDhCAiNB7ASt43WOlEcyvb2pNu4IODoPz9/n1xv9vOkForiPDdwXVagiI23&#xA;jocG3rdPxR0oScgWJ8/UwIEENu3lRb09BBKV0ywNjCmFe63FjDSjlD4Q50nv44wpl3HQAc0OB9I3&#xA;VXWrjedEnIq03oC27qjHjympx31un0wiTQLq4bS4shrAECZnMGVTvbD9QSnOhDInjzsyuiBXWQP0&#xA;r4R18jfgnHvJr4nOknrFmJ5xzhraRtHWf3Kkpfz/yGo5S+H73hkquyhdjr3R8m2lj09svwAAAP//&#xA;AwBQSwMEFAAGAAgAAAAhAFr0LFu/AAAAFQEAAAsAAABfcmVscy8ucmVsc2zPwWrDMAwG4Ptg72B0&#xA;X5TuUMaI01uh19I+gLGVxCy2jGSy9e1nemrHjpL4P0nD4SetZiPRyNnCruvBUPYcYp4tXC/Htw8w&#xA;Wl0ObuVMFm6kcBhfX4Yzra62kC6xqGlKVgtLreUTUf1CyWnHhXKbTCzJ1VbKjMX5LzcTvvf9HuXR&#xA;gPHJNKdgQU5hB+ZyK23zHztFL6w81c5zQp6m6P9TMfB3PtPWFCczVQtB9N4U2rp2HOA44NMz4y8A&#xA;AAD//wMAUEsDBBQABgAIAAAAIQD5LUR9wgAAAN0AAAAPAAAAZHJzL2Rvd25yZXYueG1sRI9Bi8Iw&#xA;FITvgv8hvAVvmq6K61ajFEHwanVhj4/m2Xa3ealJ1PrvjSB4HGbmG2a57kwjruR8bVnB5ygBQVxY&#xA;XXOp4HjYDucgfEDW2FgmBXfysF71e0tMtb3xnq55KEWEsE9RQRVCm0rpi4oM+pFtiaN3ss5giNKV&#xA;Uju8Rbhp5DhJZtJgzXGhwpY2FRX/+cUoyLK/7uecf+PWy3niZnqqy+xXqcFHly1ABOrCO/xq77SC&#xA;8dd0As838QnI1QMAAP//AwBQSwECLQAUAAYACAAAACEA2+H2y+4AAACFAQAAEwAAAAAAAAAAAAAA&#xA;AAAAAAAAW0NvbnRlbnRfVHlwZXNdLnhtbFBLAQItABQABgAIAAAAIQBa9CxbvwAAABUBAAALAAAA&#xA;AAAAAAAAAAAAAB8BAABfcmVscy8ucmVsc1BLAQItABQABgAIAAAAIQD5LUR9wgAAAN0AAAAPAAAA&#xA;AAAAAAAAAAAAAAcCAABkcnMvZG93bnJldi54bWxQSwUGAAAAAAMAAwC3AAAA9gIAAAAA&#xA;" filled="f" stroked="f" strokeweight=".25pt"><v:textbox style="mso-next-textbox:#Rectangle 954" inset="1pt,1pt,1pt,1pt"><w:txbxContent><w:p w14:paraId="649F74BE" w14:textId="77777777" w:rsidR="00B827FD" w:rsidRPr="00951A0B" w:rsidRDefault="00B827FD" w:rsidP="00B827FD"><w:pPr><w:pStyle w:val="a1"/><w:jc w:val="left"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr><w:r w:rsidRPr="00951A0B"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>Н. Контр.</w:t></w:r></w:p></w:txbxContent></v:textbox></v:rect><v:rect id="Rectangle 955" o:spid="_x0000_s1996" style="position:absolute;left:9281;width:10718;height:20000;visibility:visible" o:gfxdata="UEsDBBQABgAIAAAAIQDb4fbL7gAAAIUBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbHyQz07DMAyH&#xA;70i8Q+QralM4IITa7kDhCAiNB7ASt43WOlEcyvb2pNu4IODoPz9/n1xv9vOkForiPDdwXVagiI23&#xA;jocG3rdPxR0oScgWJ8/UwIEENu3lRb09BBKV0ywNjCmFe63FjDSjlD4Q50nv44wpl3HQAc0OB9I3&#xA;VXWrjedEnIq03oC27qjHjympx31un0wiTQLq4bS4shrAECZnMGVTvbD9QSnOhDInjzsyuiBXWQP0&#xA;r4R18jfgnHvJr4nOknrFmJ5xzhraRtHWf3Kkpfz/yGo5S+H73hkquyhdjr3R8m2lj09svwAAAP//&#xA;AwBQSwMEFAAGAAgAAAAhAFr0LFu/AAAAFQEAAAsAAABfcmVscy8ucmVsc2zPwWrDMAwG4Ptg72B0&#xA;X5TuUMaI01uh19I+gLGVxCy2jGSy9e1nemrHjpL4P0nD4SetZiPRyNnCruvBUPYcYp4tXC/Htw8w&#xA;Wl0ObuVMFm6kcBhfX4Yzra62kC6xqGlKVgtLreUTUf1CyWnHhXKbTCzJ1VbKjMX5LzcTvvf9HuXR&#xA;gPHJNKdgQU5hB+ZyK23zHztFL6w81c5zQp6m6P9TMfB3PtPWFCczVQtB9N4U2rp2HOA44NMz4y8A&#xA;AAD//wMAUEsDBBQABgAIAAAAIQB2xNwJwgAAAN0AAAAPAAAAZHJzL2Rvd25yZXYueG1sRI9Bi8Iw&#xA;FITvC/6H8ARva6oUV6tRiiDs1boLe3w0z7bavNQkq/XfG0HwOMzMN8xq05tWXMn5xrKCyTgBQVxa&#xA;3XCl4Oew+5yD8AFZY2uZFNzJw2Y9+Fhhpu2N93QtQiUihH2GCuoQukxKX9Zk0I9tRxy9o3UGQ5Su&#xA;ktrhLcJNK6dJMpMGG44LNXa0rak8F/9GQZ6f+t9LscCdl/PEzXSqq/xPqdGwz5cgAvXhHX61v7WC&#xA;6VeawvNNfAJy/QAAAP//AwBQSwECLQAUAAYACAAAACEA2+H2y+4AAACFAQAAEwAAAAAAAAAAAAAA&#xA;AAAAAAAAW0NvbnRlbnRfVHlwZXNdLnhtbFBLAQItABQABgAIAAAAIQBa9CxbvwAAABUBAAALAAAA&#xA;AAAAAAAAAAAAAB8BAABfcmVscy8ucmVsc1BLAQItABQABgAIAAAAIQB2xNwJwgAAAN0AAAAPAAAA&#xA;AAAAAAAAAAAAAAcCAABkcnMvZG93bnJldi54bWxQSwUGAAAAAAMAAwC3AAAA9gIAAAAA&#xA;" filled="f" stroked="f" strokeweight=".25pt"><v:textbox style="mso-next-textbox:#Rectangle 955" inset="1pt,1pt,1pt,1pt"><w:txbxContent><w:p w14:paraId="544E785C" w14:textId="77777777" w:rsidR="00B827FD" w:rsidRPr="00951A0B" w:rsidRDefault="00B827FD" w:rsidP="00B827FD"><w:pPr><w:rPr><w:i/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="ru-RU"/></w:rPr></w:pPr><w:r w:rsidRPr="00951A0B"><w:rPr><w:i/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="ru-RU"/></w:rPr><w:t>Рыжанкова А.С.</w:t></w:r></w:p></w:txbxContent></v:textbox></v:rect></v:group><v:group id="Group 956" o:spid="_x0000_s1997" style="position:absolute;left:1144;top:16014;width:2507;height:240" coordsize="19999,20000" o:gfxdata="UEsDBBQABgAIAAAAIQDb4fbL7gAAAIUBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbHyQz07DMAyH&#xA;70i8Q+QralM4IITa7kDhCAiNB7ASt43WOlEcyvb2pNu4IODoPz9/n1xv9vOkForiPDdwXVagiI23&#xA;jocG3rdPxR0oScgWJ8/UwIEENu3lRb09BBKV0ywNjCmFe63FjDSjlD4Q50nv44wpl3HQAc0OB9I3&#xA;VXWrjedEnIq03oC27qjHjympx31un0wiTQLq4bS4shrAECZnMGVTvbD9QSnOhDInjzsyuiBXWQP0&#xA;r4R18jfgnHvJr4nOknrFmJ5xzhraRtHWf3Kkpfz/yGo5S+H73hkquyhdjr3R8m2lj09svwAAAP//&#xA;AwBQSwMEFAAGAAgAAAAhAFr0LFu/AAAAFQEAAAsAAABfcmVscy8ucmVsc2zPwWrDMAwG4Ptg72B0&#xA;X5TuUMaI01uh19I+gLGVxCy2jGSy9e1nemrHjpL4P0nD4SetZiPRyNnCruvBUPYcYp4tXC/Htw8w&#xA;Wl0ObuVMFm6kcBhfX4Yzra62kC6xqGlKVgtLreUTUf1CyWnHhXKbTCzJ1VbKjMX5LzcTvvf9HuXR&#xA;gPHJNKdgQU5hB+ZyK23zHztFL6w81c5zQp6m6P9TMfB3PtPWFCczVQtB9N4U2rp2HOA44NMz4y8A&#xA;AAD//wMAUEsDBBQABgAIAAAAIQBQR8NPxgAAAN0AAAAPAAAAZHJzL2Rvd25yZXYueG1sRI9Ba8JA&#xA;FITvhf6H5RW86SZaW4muIqLiQYRqQbw9ss8kmH0bsmsS/31XEHocZuYbZrboTCkaql1hWUE8iEAQ&#xA;p1YXnCn4PW36ExDOI2ssLZOCBzlYzN/fZpho2/IPNUefiQBhl6CC3PsqkdKlORl0A1sRB+9qa4M+&#xA;yDqTusY2wE0ph1H0JQ0WHBZyrGiVU3o73o2CbYvtchSvm/3tunpcTuPDeR+TUr2PbjkF4anz/+FX&#xA;e6cVDL8/x/B8E56AnP8BAAD//wMAUEsBAi0AFAAGAAgAAAAhANvh9svuAAAAhQEAABMAAAAAAAAA&#xA;AAAAAAAAAAAAAFtDb250ZW50X1R5cGVzXS54bWxQSwECLQAUAAYACAAAACEAWvQsW78AAAAVAQAA&#xA;CwAAAAAAAAAAAAAAAAAfAQAAX3JlbHMvLnJlbHNQSwECLQAUAAYACAAAACEAUEfDT8YAAADdAAAA&#xA;DwAAAAAAAAAAAAAAAAAHAgAAZHJzL2Rvd25yZXYueG1sUEsFBgAAAAADAAMAtwAAAPoCAAAAAA==&#xA;"><v:rect id="Rectangle 957" o:spid="_x0000_s1998" style="position:absolute;width:8856;height:20000;visibility:visible" o:gfxdata="UEsDBBQABgAIAAAAIQDb4fbL7gAAAIUBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbHyQz07DMAyH&#xA;70i8Q+QralM4IITa7kDhCAiNB7ASt43WOlEcyvb2pNu4IODoPz9/n1xv9vOkForiPDdwXVagiI23&#xA;jocG3rdPxR0oScgWJ8/UwIEENu3lRb09BBKV0ywNjCmFe63FjDSjlD4Q50nv44wpl3HQAc0OB9I3&#xA;VXWrjedEnIq03oC27qjHjympx31un0wiTQLq4bS4shrAECZnMGVTvbD9QSnOhDInjzsyuiBXWQP0&#xA;r4R18jfgnHvJr4nOknrFmJ5xzhraRtHWf3Kkpfz/yGo5S+H73hkquyhdjr3R8m2lj09svwAAAP//&#xA;AwBQSwMEFAAGAAgAAAAhAFr0LFu/AAAAFQEAAAsAAABfcmVscy8ucmVsc2zPwWrDMAwG4Ptg72B0&#xA;X5TuUMaI01uh19I+gLGVxCy2jGSy9e1nemrHjpL4P0nD4SetZiPRyNnCruvBUPYcYp4tXC/Htw8w&#xA;Wl0ObuVMFm6kcBhfX4Yzra62kC6xqGlKVgtLreUTUf1CyWnHhXKbTCzJ1VbKjMX5LzcTvvf9HuXR&#xA;gPHJNKdgQU5hB+ZyK23zHztFL6w81c5zQp6m6P9TMfB3PtPWFCczVQtB9N4U2rp2HOA44NMz4y8A&#xA;AAD//wMAUEsDBBQABgAIAAAAIQDpWuflwgAAAN0AAAAPAAAAZHJzL2Rvd25yZXYueG1sRI9Bi8Iw&#xA;FITvC/6H8ARva6pIV6tRiiDs1boLe3w0z7bavNQkq/XfG0HwOMzMN8xq05tWXMn5xrKCyTgBQVxa&#xA;3XCl4Oew+5yD8AFZY2uZFNzJw2Y9+Fhhpu2N93QtQiUihH2GCuoQukxKX9Zk0I9tRxy9o3UGQ5Su&#xA;ktrhLcJNK6dJkkqDDceFGjva1lSei3+jIM9P/e+lWODOy3niUj3TVf6n1GjY50sQgfrwDr/a31rB&#xA;9GuWwvNNfAJy/QAAAP//AwBQSwECLQAUAAYACAAAACEA2+H2y+4AAACFAQAAEwAAAAAAAAAAAAAA&#xA;AAAAAAAAW0NvbnRlbnRfVHlwZXNdLnhtbFBLAQItABQABgAIAAAAIQBa9CxbvwAAABUBAAALAAAA&#xA;AAAAAAAAAAAAAB8BAABfcmVscy8ucmVsc1BLAQItABQABgAIAAAAIQDpWuflwgAAAN0AAAAPAAAA&#xA;AAAAAAAAAAAAAAcCAABkcnMvZG93bnJldi54bWxQSwUGAAAAAAMAAwC3AAAA9gIAAAAA&#xA;" filled="f" stroked="f" strokeweight=".25pt"><v:textbox style="mso-next-textbox:#Rectangle 957" inset="1pt,1pt,1pt,1pt"><w:txbxContent><w:p w14:paraId="616C4B5E" w14:textId="77777777" w:rsidR="00B827FD" w:rsidRPr="00951A0B" w:rsidRDefault="00B827FD" w:rsidP="00B827FD"><w:pPr><w:pStyle w:val="a1"/><w:jc w:val="left"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr><w:r w:rsidRPr="00951A0B"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>Утверд.</w:t></w:r></w:p></w:txbxContent></v:textbox></v:rect><v:rect id="Rectangle 958" o:spid="_x0000_s1999" style="position:absolute;left:9281;width:10718;height:20000;visibility:visible" o:gfxdata="UEsDBBQABgAIAAAAIQDb4fbL7gAAAIUBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbHyQz07DMAyH&#xA;70i8Q+QralM4IITa7kDhCAiNB7ASt43WOlEcyvb2pNu4IODoPz9/n1xv9vOkForiPDdwXVagiI23&#xA;jocG3rdPxR0oScgWJ8/UwIEENu3lRb09BBKV0ywNjCmFe63FjDSjlD4Q50nv44wpl3HQAc0OB9I3&#xA;VXWrjedEnIq03oC27qjHjympx31un0wiTQLq4bS4shrAECZnMGVTvbD9QSnOhDInjzsyuiBXWQP0&#xA;r4R18jfgnHvJr4nOknrFmJ5xzhraRtHWf3Kkpfz/yGo5S+H73hkquyhdjr3R8m2lj09svwAAAP//&#xA;AwBQSwMEFAAGAAgAAAAhAFr0LFu/AAAAFQEAAAsAAABfcmVscy8ucmVsc2zPwWrDMAwG4Ptg72B0&#xA;X5TuUMaI01uh19I+gLGVxCy2jGSy9e1nemrHjpL4P0nD4SetZiPRyNnCruvBUPYcYp4tXC/Htw8w&#xA;Wl0ObuVMFm6kcBhfX4Yzra62kC6xqGlKVgtLreUTUf1CyWnHhXKbTCzJ1VbKjMX5LzcTvvf9HuXR&#xA;gPHJNKdgQU5hB+ZyK23zHztFL6w81c5zQp6m6P9TMfB3PtPWFCczVQtB9N4U2rp2HOA44NMz4y8A&#xA;AAD//wMAUEsDBBQABgAIAAAAIQCGFkJ+wgAAAN0AAAAPAAAAZHJzL2Rvd25yZXYueG1sRI9Bi8Iw&#xA;FITvwv6H8Ba8aboiWrtGKYLgdauCx0fztq02L90kavffG0HwOMzMN8xy3ZtW3Mj5xrKCr3ECgri0&#xA;uuFKwWG/HaUgfEDW2FomBf/kYb36GCwx0/bOP3QrQiUihH2GCuoQukxKX9Zk0I9tRxy9X+sMhihd&#xA;JbXDe4SbVk6SZCYNNhwXauxoU1N5Ka5GQZ6f++NfscCtl2niZnqqq/yk1PCzz79BBOrDO/xq77SC&#xA;yXw6h+eb+ATk6gEAAP//AwBQSwECLQAUAAYACAAAACEA2+H2y+4AAACFAQAAEwAAAAAAAAAAAAAA&#xA;AAAAAAAAW0NvbnRlbnRfVHlwZXNdLnhtbFBLAQItABQABgAIAAAAIQBa9CxbvwAAABUBAAALAAAA&#xA;AAAAAAAAAAAAAB8BAABfcmVscy8ucmVsc1BLAQItABQABgAIAAAAIQCGFkJ+wgAAAN0AAAAPAAAA&#xA;AAAAAAAAAAAAAAcCAABkcnMvZG93bnJldi54bWxQSwUGAAAAAAMAAwC3AAAA9gIAAAAA&#xA;" filled="f" stroked="f" strokeweight=".25pt"><v:textbox style="mso-next-textbox:#Rectangle 958" inset="1pt,1pt,1pt,1pt"><w:txbxContent><w:p w14:paraId="12944F36" w14:textId="77777777" w:rsidR="00B827FD" w:rsidRPr="00CF1877" w:rsidRDefault="00B827FD" w:rsidP="00B827FD"><w:pPr><w:rPr><w:i/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="ru-RU"/></w:rPr></w:pPr><w:r><w:rPr><w:i/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="ru-RU"/></w:rPr><w:t>Пацей Н.В.</w:t></w:r></w:p></w:txbxContent></v:textbox></v:rect></v:group><v:line id="Line 959" o:spid="_x0000_s2000" style="position:absolute;visibility:visible" from="8541,14913" to="8542,16262" o:connectortype="straight" o:gfxdata="UEsDBBQABgAIAAAAIQDb4fbL7gAAAIUBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbHyQz07DMAyH&#xA;70i8Q+QralM4IITa7kDhCAiNB7ASt43WOlEcyvb2pNu4IODoPz9/n1xv9vOkForiPDdwXVagiI23&#xA;jocG3rdPxR0oScgWJ8/UwIEENu3lRb09BBKV0ywNjCmFe63FjDSjlD4Q50nv44wpl3HQAc0OB9I3&#xA;VXWrjedEnIq03oC27qjHjympx31un0wiTQLq4bS4shrAECZnMGVTvbD9QSnOhDInjzsyuiBXWQP0&#xA;r4R18jfgnHvJr4nOknrFmJ5xzhraRtHWf3Kkpfz/yGo5S+H73hkquyhdjr3R8m2lj09svwAAAP//&#xA;AwBQSwMEFAAGAAgAAAAhAFr0LFu/AAAAFQEAAAsAAABfcmVscy8ucmVsc2zPwWrDMAwG4Ptg72B0&#xA;X5TuUMaI01uh19I+gLGVxCy2jGSy9e1nemrHjpL4P0nD4SetZiPRyNnCruvBUPYcYp4tXC/Htw8w&#xA;Wl0ObuVMFm6kcBhfX4Yzra62kC6xqGlKVgtLreUTUf1CyWnHhXKbTCzJ1VbKjMX5LzcTvvf9HuXR&#xA;gPHJNKdgQU5hB+ZyK23zHztFL6w81c5zQp6m6P9TMfB3PtPWFCczVQtB9N4U2rp2HOA44NMz4y8A&#xA;AAD//wMAUEsDBBQABgAIAAAAIQBQxTV6vgAAAN0AAAAPAAAAZHJzL2Rvd25yZXYueG1sRE+7CsIw&#xA;FN0F/yFcwU1TxRfVKCJU3MTq4nZtrm2xuSlN1Pr3ZhAcD+e92rSmEi9qXGlZwWgYgSDOrC45V3A5&#xA;J4MFCOeRNVaWScGHHGzW3c4KY23ffKJX6nMRQtjFqKDwvo6ldFlBBt3Q1sSBu9vGoA+wyaVu8B3C&#xA;TSXHUTSTBksODQXWtCsoe6RPo+BxvUyT/XGnz1W61bc88dfbXSvV77XbJQhPrf+Lf+6DVjCeT8Lc&#xA;8CY8Abn+AgAA//8DAFBLAQItABQABgAIAAAAIQDb4fbL7gAAAIUBAAATAAAAAAAAAAAAAAAAAAAA&#xA;AABbQ29udGVudF9UeXBlc10ueG1sUEsBAi0AFAAGAAgAAAAhAFr0LFu/AAAAFQEAAAsAAAAAAAAA&#xA;AAAAAAAAHwEAAF9yZWxzLy5yZWxzUEsBAi0AFAAGAAgAAAAhAFDFNXq+AAAA3QAAAA8AAAAAAAAA&#xA;AAAAAAAABwIAAGRycy9kb3ducmV2LnhtbFBLBQYAAAAAAwADALcAAADyAgAAAAA=&#xA;" strokeweight="2pt"/><v:rect id="Rectangle 960" o:spid="_x0000_s2001" style="position:absolute;left:5189;top:15211;width:3285;height:1007;visibility:visible" o:gfxdata="UEsDBBQABgAIAAAAIQDb4fbL7gAAAIUBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbHyQz07DMAyH&#xA;70i8Q+QralM4IITa7kDhCAiNB7ASt43WOlEcyvb2pNu4IODoPz9/n1xv9vOkForiPDdwXVagiI23&#xA;jocG3rdPxR0oScgWJ8/UwIEENu3lRb09BBKV0ywNjCmFe63FjDSjlD4Q50nv44wpl3HQAc0OB9I3&#xA;VXWrjedEnIq03oC27qjHjympx31un0wiTQLq4bS4shrAECZnMGVTvbD9QSnOhDInjzsyuiBXWQP0&#xA;r4R18jfgnHvJr4nOknrFmJ5xzhraRtHWf3Kkpfz/yGo5S+H73hkquyhdjr3R8m2lj09svwAAAP//&#xA;AwBQSwMEFAAGAAgAAAAhAFr0LFu/AAAAFQEAAAsAAABfcmVscy8ucmVsc2zPwWrDMAwG4Ptg72B0&#xA;X5TuUMaI01uh19I+gLGVxCy2jGSy9e1nemrHjpL4P0nD4SetZiPRyNnCruvBUPYcYp4tXC/Htw8w&#xA;Wl0ObuVMFm6kcBhfX4Yzra62kC6xqGlKVgtLreUTUf1CyWnHhXKbTCzJ1VbKjMX5LzcTvvf9HuXR&#xA;gPHJNKdgQU5hB+ZyK23zHztFL6w81c5zQp6m6P9TMfB3PtPWFCczVQtB9N4U2rp2HOA44NMz4y8A&#xA;AAD//wMAUEsDBBQABgAIAAAAIQCYxXOXxAAAAN0AAAAPAAAAZHJzL2Rvd25yZXYueG1sRI9Ba8JA&#xA;FITvBf/D8oTe6kYJNomuEgqCV9MKPT6yr0k0+zbdXU36791CocdhZr5htvvJ9OJOzneWFSwXCQji&#xA;2uqOGwUf74eXDIQPyBp7y6Tghzzsd7OnLRbajnyiexUaESHsC1TQhjAUUvq6JYN+YQfi6H1ZZzBE&#xA;6RqpHY4Rbnq5SpK1NNhxXGhxoLeW6mt1MwrK8jKdv6scD15miVvrVDflp1LP86ncgAg0hf/wX/uo&#xA;Faxe0xx+38QnIHcPAAAA//8DAFBLAQItABQABgAIAAAAIQDb4fbL7gAAAIUBAAATAAAAAAAAAAAA&#xA;AAAAAAAAAABbQ29udGVudF9UeXBlc10ueG1sUEsBAi0AFAAGAAgAAAAhAFr0LFu/AAAAFQEAAAsA&#xA;AAAAAAAAAAAAAAAAHwEAAF9yZWxzLy5yZWxzUEsBAi0AFAAGAAgAAAAhAJjFc5fEAAAA3QAAAA8A&#xA;AAAAAAAAAAAAAAAABwIAAGRycy9kb3ducmV2LnhtbFBLBQYAAAAAAwADALcAAAD4AgAAAAA=&#xA;" filled="f" stroked="f" strokeweight=".25pt"><v:textbox style="mso-next-textbox:#Rectangle 960" inset="1pt,1pt,1pt,1pt"><w:txbxContent><w:p w14:paraId="5EAF79F8" w14:textId="77777777" w:rsidR="00B827FD" w:rsidRPr="00CB455E" w:rsidRDefault="00B827FD" w:rsidP="00B827FD"><w:pPr><w:pStyle w:val="a1"/><w:suppressAutoHyphens/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="ru-RU"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="ru-RU"/></w:rPr><w:t>Скриншоты работы программы</w:t></w:r></w:p></w:txbxContent></v:textbox></v:rect><v:line id="Line 961" o:spid="_x0000_s2002" style="position:absolute;visibility:visible" from="8548,15183" to="11560,15183" o:connectortype="straight" o:gfxdata="UEsDBBQABgAIAAAAIQDb4fbL7gAAAIUBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbHyQz07DMAyH&#xA;70i8Q+QralM4IITa7kDhCAiNB7ASt43WOlEcyvb2pNu4IODoPz9/n1xv9vOkForiPDdwXVagiI23&#xA;jocG3rdPxR0oScgWJ8/UwIEENu3lRb09BBKV0ywNjCmFe63FjDSjlD4Q50nv44wpl3HQAc0OB9I3&#xA;VXWrjedEnIq03oC27qjHjympx31un0wiTQLq4bS4shrAECZnMGVTvbD9QSnOhDInjzsyuiBXWQP0&#xA;r4R18jfgnHvJr4nOknrFmJ5xzhraRtHWf3Kkpfz/yGo5S+H73hkquyhdjr3R8m2lj09svwAAAP//&#xA;AwBQSwMEFAAGAAgAAAAhAFr0LFu/AAAAFQEAAAsAAABfcmVscy8ucmVsc2zPwWrDMAwG4Ptg72B0&#xA;X5TuUMaI01uh19I+gLGVxCy2jGSy9e1nemrHjpL4P0nD4SetZiPRyNnCruvBUPYcYp4tXC/Htw8w&#xA;Wl0ObuVMFm6kcBhfX4Yzra62kC6xqGlKVgtLreUTUf1CyWnHhXKbTCzJ1VbKjMX5LzcTvvf9HuXR&#xA;gPHJNKdgQU5hB+ZyK23zHztFL6w81c5zQp6m6P9TMfB3PtPWFCczVQtB9N4U2rp2HOA44NMz4y8A&#xA;AAD//wMAUEsDBBQABgAIAAAAIQAraq+hvgAAAN0AAAAPAAAAZHJzL2Rvd25yZXYueG1sRE+7CsIw&#xA;FN0F/yFcwU1TBR9Uo4hQcRNrF7drc22LzU1pota/N4PgeDjv9bYztXhR6yrLCibjCARxbnXFhYLs&#xA;koyWIJxH1lhbJgUfcrDd9HtrjLV985leqS9ECGEXo4LS+yaW0uUlGXRj2xAH7m5bgz7AtpC6xXcI&#xA;N7WcRtFcGqw4NJTY0L6k/JE+jYLHNZslh9NeX+p0p29F4q+3u1ZqOOh2KxCeOv8X/9xHrWC6mIX9&#xA;4U14AnLzBQAA//8DAFBLAQItABQABgAIAAAAIQDb4fbL7gAAAIUBAAATAAAAAAAAAAAAAAAAAAAA&#xA;AABbQ29udGVudF9UeXBlc10ueG1sUEsBAi0AFAAGAAgAAAAhAFr0LFu/AAAAFQEAAAsAAAAAAAAA&#xA;AAAAAAAAHwEAAF9yZWxzLy5yZWxzUEsBAi0AFAAGAAgAAAAhACtqr6G+AAAA3QAAAA8AAAAAAAAA&#xA;AAAAAAAABwIAAGRycy9kb3ducmV2LnhtbFBLBQYAAAAAAwADALcAAADyAgAAAAA=&#xA;" strokeweight="2pt"/><v:line id="Line 962" o:spid="_x0000_s2003" style="position:absolute;visibility:visible" from="8547,15456" to="11559,15457" o:connectortype="straight" o:gfxdata="UEsDBBQABgAIAAAAIQDb4fbL7gAAAIUBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbHyQz07DMAyH&#xA;70i8Q+QralM4IITa7kDhCAiNB7ASt43WOlEcyvb2pNu4IODoPz9/n1xv9vOkForiPDdwXVagiI23&#xA;jocG3rdPxR0oScgWJ8/UwIEENu3lRb09BBKV0ywNjCmFe63FjDSjlD4Q50nv44wpl3HQAc0OB9I3&#xA;VXWrjedEnIq03oC27qjHjympx31un0wiTQLq4bS4shrAECZnMGVTvbD9QSnOhDInjzsyuiBXWQP0&#xA;r4R18jfgnHvJr4nOknrFmJ5xzhraRtHWf3Kkpfz/yGo5S+H73hkquyhdjr3R8m2lj09svwAAAP//&#xA;AwBQSwMEFAAGAAgAAAAhAFr0LFu/AAAAFQEAAAsAAABfcmVscy8ucmVsc2zPwWrDMAwG4Ptg72B0&#xA;X5TuUMaI01uh19I+gLGVxCy2jGSy9e1nemrHjpL4P0nD4SetZiPRyNnCruvBUPYcYp4tXC/Htw8w&#xA;Wl0ObuVMFm6kcBhfX4Yzra62kC6xqGlKVgtLreUTUf1CyWnHhXKbTCzJ1VbKjMX5LzcTvvf9HuXR&#xA;gPHJNKdgQU5hB+ZyK23zHztFL6w81c5zQp6m6P9TMfB3PtPWFCczVQtB9N4U2rp2HOA44NMz4y8A&#xA;AAD//wMAUEsDBBQABgAIAAAAIQBEJgo6wQAAAN0AAAAPAAAAZHJzL2Rvd25yZXYueG1sRI/NCsIw&#xA;EITvgu8QVvCmqYI/VKOIUPEmVi/e1mZti82mNFHr2xtB8DjMzDfMct2aSjypcaVlBaNhBII4s7rk&#xA;XMH5lAzmIJxH1lhZJgVvcrBedTtLjLV98ZGeqc9FgLCLUUHhfR1L6bKCDLqhrYmDd7ONQR9kk0vd&#xA;4CvATSXHUTSVBksOCwXWtC0ou6cPo+B+OU+S3WGrT1W60dc88ZfrTSvV77WbBQhPrf+Hf+29VjCe&#xA;TUbwfROegFx9AAAA//8DAFBLAQItABQABgAIAAAAIQDb4fbL7gAAAIUBAAATAAAAAAAAAAAAAAAA&#xA;AAAAAABbQ29udGVudF9UeXBlc10ueG1sUEsBAi0AFAAGAAgAAAAhAFr0LFu/AAAAFQEAAAsAAAAA&#xA;AAAAAAAAAAAAHwEAAF9yZWxzLy5yZWxzUEsBAi0AFAAGAAgAAAAhAEQmCjrBAAAA3QAAAA8AAAAA&#xA;AAAAAAAAAAAABwIAAGRycy9kb3ducmV2LnhtbFBLBQYAAAAAAwADALcAAAD1AgAAAAA=&#xA;" strokeweight="2pt"/><v:line id="Line 963" o:spid="_x0000_s2004" style="position:absolute;visibility:visible" from="10253,14913" to="10255,15450" o:connectortype="straight" o:gfxdata="UEsDBBQABgAIAAAAIQDb4fbL7gAAAIUBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbHyQz07DMAyH&#xA;70i8Q+QralM4IITa7kDhCAiNB7ASt43WOlEcyvb2pNu4IODoPz9/n1xv9vOkForiPDdwXVagiI23&#xA;jocG3rdPxR0oScgWJ8/UwIEENu3lRb09BBKV0ywNjCmFe63FjDSjlD4Q50nv44wpl3HQAc0OB9I3&#xA;VXWrjedEnIq03oC27qjHjympx31un0wiTQLq4bS4shrAECZnMGVTvbD9QSnOhDInjzsyuiBXWQP0&#xA;r4R18jfgnHvJr4nOknrFmJ5xzhraRtHWf3Kkpfz/yGo5S+H73hkquyhdjr3R8m2lj09svwAAAP//&#xA;AwBQSwMEFAAGAAgAAAAhAFr0LFu/AAAAFQEAAAsAAABfcmVscy8ucmVsc2zPwWrDMAwG4Ptg72B0&#xA;X5TuUMaI01uh19I+gLGVxCy2jGSy9e1nemrHjpL4P0nD4SetZiPRyNnCruvBUPYcYp4tXC/Htw8w&#xA;Wl0ObuVMFm6kcBhfX4Yzra62kC6xqGlKVgtLreUTUf1CyWnHhXKbTCzJ1VbKjMX5LzcTvvf9HuXR&#xA;gPHJNKdgQU5hB+ZyK23zHztFL6w81c5zQp6m6P9TMfB3PtPWFCczVQtB9N4U2rp2HOA44NMz4y8A&#xA;AAD//wMAUEsDBBQABgAIAAAAIQC09JRNxAAAAN0AAAAPAAAAZHJzL2Rvd25yZXYueG1sRI9Bi8Iw&#xA;FITvwv6H8IS9aWpBXbqNIkLF22LtxduzebalzUtpstr99xtB8DjMzDdMuh1NJ+40uMaygsU8AkFc&#xA;Wt1wpaA4Z7MvEM4ja+wsk4I/crDdfExSTLR98Inuua9EgLBLUEHtfZ9I6cqaDLq57YmDd7ODQR/k&#xA;UEk94CPATSfjKFpJgw2HhRp72tdUtvmvUdBeimV2+Nnrc5fv9LXK/OV600p9TsfdNwhPo3+HX+2j&#xA;VhCvlzE834QnIDf/AAAA//8DAFBLAQItABQABgAIAAAAIQDb4fbL7gAAAIUBAAATAAAAAAAAAAAA&#xA;AAAAAAAAAABbQ29udGVudF9UeXBlc10ueG1sUEsBAi0AFAAGAAgAAAAhAFr0LFu/AAAAFQEAAAsA&#xA;AAAAAAAAAAAAAAAAHwEAAF9yZWxzLy5yZWxzUEsBAi0AFAAGAAgAAAAhALT0lE3EAAAA3QAAAA8A&#xA;AAAAAAAAAAAAAAAABwIAAGRycy9kb3ducmV2LnhtbFBLBQYAAAAAAwADALcAAAD4AgAAAAA=&#xA;" strokeweight="2pt"/><v:rect id="Rectangle 964" o:spid="_x0000_s2005" style="position:absolute;left:8587;top:14928;width:769;height:239;visibility:visible" o:gfxdata="UEsDBBQABgAIAAAAIQDb4fbL7gAAAIUBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbHyQz07DMAyH&#xA;70i8Q+QralM4IITa7kDhCAiNB7ASt43WOlEcyvb2pNu4IODoPz9/n1xv9vOkForiPDdwXVagiI23&#xA;jocG3rdPxR0oScgWJ8/UwIEENu3lRb09BBKV0ywNjCmFe63FjDSjlD4Q50nv44wpl3HQAc0OB9I3&#xA;VXWrjedEnIq03oC27qjHjympx31un0wiTQLq4bS4shrAECZnMGVTvbD9QSnOhDInjzsyuiBXWQP0&#xA;r4R18jfgnHvJr4nOknrFmJ5xzhraRtHWf3Kkpfz/yGo5S+H73hkquyhdjr3R8m2lj09svwAAAP//&#xA;AwBQSwMEFAAGAAgAAAAhAFr0LFu/AAAAFQEAAAsAAABfcmVscy8ucmVsc2zPwWrDMAwG4Ptg72B0&#xA;X5TuUMaI01uh19I+gLGVxCy2jGSy9e1nemrHjpL4P0nD4SetZiPRyNnCruvBUPYcYp4tXC/Htw8w&#xA;Wl0ObuVMFm6kcBhfX4Yzra62kC6xqGlKVgtLreUTUf1CyWnHhXKbTCzJ1VbKjMX5LzcTvvf9HuXR&#xA;gPHJNKdgQU5hB+ZyK23zHztFL6w81c5zQp6m6P9TMfB3PtPWFCczVQtB9N4U2rp2HOA44NMz4y8A&#xA;AAD//wMAUEsDBBQABgAIAAAAIQB89NKgxAAAAN0AAAAPAAAAZHJzL2Rvd25yZXYueG1sRI9Ba8JA&#xA;FITvgv9heUJvdaNtU02zkVAQem1U6PGRfSap2bdxd9X033cLBY/DzHzD5JvR9OJKzneWFSzmCQji&#xA;2uqOGwX73fZxBcIHZI29ZVLwQx42xXSSY6btjT/pWoVGRAj7DBW0IQyZlL5uyaCf24E4ekfrDIYo&#xA;XSO1w1uEm14ukySVBjuOCy0O9N5SfaouRkFZfo+Hc7XGrZerxKX6WTfll1IPs7F8AxFoDPfwf/tD&#xA;K1i+vjzB35v4BGTxCwAA//8DAFBLAQItABQABgAIAAAAIQDb4fbL7gAAAIUBAAATAAAAAAAAAAAA&#xA;AAAAAAAAAABbQ29udGVudF9UeXBlc10ueG1sUEsBAi0AFAAGAAgAAAAhAFr0LFu/AAAAFQEAAAsA&#xA;AAAAAAAAAAAAAAAAHwEAAF9yZWxzLy5yZWxzUEsBAi0AFAAGAAgAAAAhAHz00qDEAAAA3QAAAA8A&#xA;AAAAAAAAAAAAAAAABwIAAGRycy9kb3ducmV2LnhtbFBLBQYAAAAAAwADALcAAAD4AgAAAAA=&#xA;" filled="f" stroked="f" strokeweight=".25pt"><v:textbox style="mso-next-textbox:#Rectangle 964" inset="1pt,1pt,1pt,1pt"><w:txbxContent><w:p w14:paraId="0CE71782" w14:textId="77777777" w:rsidR="00B827FD" w:rsidRPr="0042148F" w:rsidRDefault="00B827FD" w:rsidP="00B827FD"><w:pPr><w:pStyle w:val="a1"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr></w:pPr><w:r w:rsidRPr="00B827FD"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>Лит</w:t></w:r><w:r w:rsidRPr="0042148F"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>.</w:t></w:r></w:p></w:txbxContent></v:textbox></v:rect><v:rect id="Rectangle 965" o:spid="_x0000_s2006" style="position:absolute;left:10300;top:14928;width:1215;height:239;visibility:visible" o:gfxdata="UEsDBBQABgAIAAAAIQDb4fbL7gAAAIUBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbHyQz07DMAyH&#xA;70i8Q+QralM4IITa7kDhCAiNB7ASt43WOlEcyvb2pNu4IODoPz9/n1xv9vOkForiPDdwXVagiI23&#xA;jocG3rdPxR0oScgWJ8/UwIEENu3lRb09BBKV0ywNjCmFe63FjDSjlD4Q50nv44wpl3HQAc0OB9I3&#xA;VXWrjedEnIq03oC27qjHjympx31un0wiTQLq4bS4shrAECZnMGVTvbD9QSnOhDInjzsyuiBXWQP0&#xA;r4R18jfgnHvJr4nOknrFmJ5xzhraRtHWf3Kkpfz/yGo5S+H73hkquyhdjr3R8m2lj09svwAAAP//&#xA;AwBQSwMEFAAGAAgAAAAhAFr0LFu/AAAAFQEAAAsAAABfcmVscy8ucmVsc2zPwWrDMAwG4Ptg72B0&#xA;X5TuUMaI01uh19I+gLGVxCy2jGSy9e1nemrHjpL4P0nD4SetZiPRyNnCruvBUPYcYp4tXC/Htw8w&#xA;Wl0ObuVMFm6kcBhfX4Yzra62kC6xqGlKVgtLreUTUf1CyWnHhXKbTCzJ1VbKjMX5LzcTvvf9HuXR&#xA;gPHJNKdgQU5hB+ZyK23zHztFL6w81c5zQp6m6P9TMfB3PtPWFCczVQtB9N4U2rp2HOA44NMz4y8A&#xA;AAD//wMAUEsDBBQABgAIAAAAIQDzHUrUxAAAAN0AAAAPAAAAZHJzL2Rvd25yZXYueG1sRI/NasMw&#xA;EITvhb6D2EBvjZzg5seJEkzB0GvdBHJcrI3txFq5kmq7b18VCj0OM/MNsz9OphMDOd9aVrCYJyCI&#xA;K6tbrhWcPornDQgfkDV2lknBN3k4Hh4f9phpO/I7DWWoRYSwz1BBE0KfSemrhgz6ue2Jo3e1zmCI&#xA;0tVSOxwj3HRymSQrabDluNBgT68NVffyyyjI89t0/iy3WHi5SdxKp7rOL0o9zaZ8ByLQFP7Df+03&#xA;rWC5fknh9018AvLwAwAA//8DAFBLAQItABQABgAIAAAAIQDb4fbL7gAAAIUBAAATAAAAAAAAAAAA&#xA;AAAAAAAAAABbQ29udGVudF9UeXBlc10ueG1sUEsBAi0AFAAGAAgAAAAhAFr0LFu/AAAAFQEAAAsA&#xA;AAAAAAAAAAAAAAAAHwEAAF9yZWxzLy5yZWxzUEsBAi0AFAAGAAgAAAAhAPMdStTEAAAA3QAAAA8A&#xA;AAAAAAAAAAAAAAAABwIAAGRycy9kb3ducmV2LnhtbFBLBQYAAAAAAwADALcAAAD4AgAAAAA=&#xA;" filled="f" stroked="f" strokeweight=".25pt"><v:textbox style="mso-next-textbox:#Rectangle 965" inset="1pt,1pt,1pt,1pt"><w:txbxContent><w:p w14:paraId="187D85BC" w14:textId="77777777" w:rsidR="00B827FD" w:rsidRPr="00B827FD" w:rsidRDefault="00B827FD" w:rsidP="00B827FD"><w:pPr><w:pStyle w:val="a1"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr><w:r w:rsidRPr="00B827FD"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>Листов</w:t></w:r></w:p></w:txbxContent></v:textbox></v:rect><v:rect id="Rectangle 966" o:spid="_x0000_s2007" style="position:absolute;left:10307;top:15203;width:1215;height:239;visibility:visible" o:gfxdata="UEsDBBQABgAIAAAAIQDb4fbL7gAAAIUBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbHyQz07DMAyH&#xA;70i8Q+QralM4IITa7kDhCAiNB7ASt43WOlEcyvb2pNu4IODoPz9/n1xv9vOkForiPDdwXVagiI23&#xA;jocG3rdPxR0oScgWJ8/UwIEENu3lRb09BBKV0ywNjCmFe63FjDSjlD4Q50nv44wpl3HQAc0OB9I3&#xA;VXWrjedEnIq03oC27qjHjympx31un0wiTQLq4bS4shrAECZnMGVTvbD9QSnOhDInjzsyuiBXWQP0&#xA;r4R18jfgnHvJr4nOknrFmJ5xzhraRtHWf3Kkpfz/yGo5S+H73hkquyhdjr3R8m2lj09svwAAAP//&#xA;AwBQSwMEFAAGAAgAAAAhAFr0LFu/AAAAFQEAAAsAAABfcmVscy8ucmVsc2zPwWrDMAwG4Ptg72B0&#xA;X5TuUMaI01uh19I+gLGVxCy2jGSy9e1nemrHjpL4P0nD4SetZiPRyNnCruvBUPYcYp4tXC/Htw8w&#xA;Wl0ObuVMFm6kcBhfX4Yzra62kC6xqGlKVgtLreUTUf1CyWnHhXKbTCzJ1VbKjMX5LzcTvvf9HuXR&#xA;gPHJNKdgQU5hB+ZyK23zHztFL6w81c5zQp6m6P9TMfB3PtPWFCczVQtB9N4U2rp2HOA44NMz4y8A&#xA;AAD//wMAUEsDBBQABgAIAAAAIQCcUe9PwgAAAN0AAAAPAAAAZHJzL2Rvd25yZXYueG1sRI9Bi8Iw&#xA;FITvgv8hvAVvmq6o61ajFEHwanVhj4/m2Xa3ealJ1PrvjSB4HGbmG2a57kwjruR8bVnB5ygBQVxY&#xA;XXOp4HjYDucgfEDW2FgmBXfysF71e0tMtb3xnq55KEWEsE9RQRVCm0rpi4oM+pFtiaN3ss5giNKV&#xA;Uju8Rbhp5DhJZtJgzXGhwpY2FRX/+cUoyLK/7uecf+PWy3niZnqiy+xXqcFHly1ABOrCO/xq77SC&#xA;8dd0Cs838QnI1QMAAP//AwBQSwECLQAUAAYACAAAACEA2+H2y+4AAACFAQAAEwAAAAAAAAAAAAAA&#xA;AAAAAAAAW0NvbnRlbnRfVHlwZXNdLnhtbFBLAQItABQABgAIAAAAIQBa9CxbvwAAABUBAAALAAAA&#xA;AAAAAAAAAAAAAB8BAABfcmVscy8ucmVsc1BLAQItABQABgAIAAAAIQCcUe9PwgAAAN0AAAAPAAAA&#xA;AAAAAAAAAAAAAAcCAABkcnMvZG93bnJldi54bWxQSwUGAAAAAAMAAwC3AAAA9gIAAAAA&#xA;" filled="f" stroked="f" strokeweight=".25pt"><v:textbox style="mso-next-textbox:#Rectangle 966" inset="1pt,1pt,1pt,1pt"><w:txbxContent><w:p w14:paraId="228036EB" w14:textId="77777777" w:rsidR="00B827FD" w:rsidRPr="00B827FD" w:rsidRDefault="00B827FD" w:rsidP="00B827FD"><w:pPr><w:pStyle w:val="a1"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="ru-RU"/></w:rPr></w:pPr><w:r w:rsidRPr="00B827FD"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="ru-RU"/></w:rPr><w:t>1</w:t></w:r></w:p></w:txbxContent></v:textbox></v:rect><v:line id="Line 967" o:spid="_x0000_s2008" style="position:absolute;visibility:visible" from="8827,15188" to="8828,15450" o:connectortype="straight" o:gfxdata="UEsDBBQABgAIAAAAIQDb4fbL7gAAAIUBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbHyQz07DMAyH&#xA;70i8Q+QralM4IITa7kDhCAiNB7ASt43WOlEcyvb2pNu4IODoPz9/n1xv9vOkForiPDdwXVagiI23&#xA;jocG3rdPxR0oScgWJ8/UwIEENu3lRb09BBKV0ywNjCmFe63FjDSjlD4Q50nv44wpl3HQAc0OB9I3&#xA;VXWrjedEnIq03oC27qjHjympx31un0wiTQLq4bS4shrAECZnMGVTvbD9QSnOhDInjzsyuiBXWQP0&#xA;r4R18jfgnHvJr4nOknrFmJ5xzhraRtHWf3Kkpfz/yGo5S+H73hkquyhdjr3R8m2lj09svwAAAP//&#xA;AwBQSwMEFAAGAAgAAAAhAFr0LFu/AAAAFQEAAAsAAABfcmVscy8ucmVsc2zPwWrDMAwG4Ptg72B0&#xA;X5TuUMaI01uh19I+gLGVxCy2jGSy9e1nemrHjpL4P0nD4SetZiPRyNnCruvBUPYcYp4tXC/Htw8w&#xA;Wl0ObuVMFm6kcBhfX4Yzra62kC6xqGlKVgtLreUTUf1CyWnHhXKbTCzJ1VbKjMX5LzcTvvf9HuXR&#xA;gPHJNKdgQU5hB+ZyK23zHztFL6w81c5zQp6m6P9TMfB3PtPWFCczVQtB9N4U2rp2HOA44NMz4y8A&#xA;AAD//wMAUEsDBBQABgAIAAAAIQBBCvSkxgAAAN0AAAAPAAAAZHJzL2Rvd25yZXYueG1sRI/dagIx&#xA;FITvC32HcAreaVZBbbdGKf6A0gvR9gGOm9PN1s3JkkRdffpGEHo5zMw3zGTW2lqcyYfKsYJ+LwNB&#xA;XDhdcang+2vVfQURIrLG2jEpuFKA2fT5aYK5dhfe0XkfS5EgHHJUYGJscilDYchi6LmGOHk/zluM&#xA;SfpSao+XBLe1HGTZSFqsOC0YbGhuqDjuT1bBxh8+j/1baeSBN35Zbxdvwf4q1XlpP95BRGrjf/jR&#xA;XmsFg/FwBPc36QnI6R8AAAD//wMAUEsBAi0AFAAGAAgAAAAhANvh9svuAAAAhQEAABMAAAAAAAAA&#xA;AAAAAAAAAAAAAFtDb250ZW50X1R5cGVzXS54bWxQSwECLQAUAAYACAAAACEAWvQsW78AAAAVAQAA&#xA;CwAAAAAAAAAAAAAAAAAfAQAAX3JlbHMvLnJlbHNQSwECLQAUAAYACAAAACEAQQr0pMYAAADdAAAA&#xA;DwAAAAAAAAAAAAAAAAAHAgAAZHJzL2Rvd25yZXYueG1sUEsFBgAAAAADAAMAtwAAAPoCAAAAAA==&#xA;" strokeweight="1pt"/><v:line id="Line 968" o:spid="_x0000_s2009" style="position:absolute;visibility:visible" from="9112,15189" to="9113,15451" o:connectortype="straight" o:gfxdata="UEsDBBQABgAIAAAAIQDb4fbL7gAAAIUBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbHyQz07DMAyH&#xA;70i8Q+QralM4IITa7kDhCAiNB7ASt43WOlEcyvb2pNu4IODoPz9/n1xv9vOkForiPDdwXVagiI23&#xA;jocG3rdPxR0oScgWJ8/UwIEENu3lRb09BBKV0ywNjCmFe63FjDSjlD4Q50nv44wpl3HQAc0OB9I3&#xA;VXWrjedEnIq03oC27qjHjympx31un0wiTQLq4bS4shrAECZnMGVTvbD9QSnOhDInjzsyuiBXWQP0&#xA;r4R18jfgnHvJr4nOknrFmJ5xzhraRtHWf3Kkpfz/yGo5S+H73hkquyhdjr3R8m2lj09svwAAAP//&#xA;AwBQSwMEFAAGAAgAAAAhAFr0LFu/AAAAFQEAAAsAAABfcmVscy8ucmVsc2zPwWrDMAwG4Ptg72B0&#xA;X5TuUMaI01uh19I+gLGVxCy2jGSy9e1nemrHjpL4P0nD4SetZiPRyNnCruvBUPYcYp4tXC/Htw8w&#xA;Wl0ObuVMFm6kcBhfX4Yzra62kC6xqGlKVgtLreUTUf1CyWnHhXKbTCzJ1VbKjMX5LzcTvvf9HuXR&#xA;gPHJNKdgQU5hB+ZyK23zHztFL6w81c5zQp6m6P9TMfB3PtPWFCczVQtB9N4U2rp2HOA44NMz4y8A&#xA;AAD//wMAUEsDBBQABgAIAAAAIQAuRlE/xQAAAN0AAAAPAAAAZHJzL2Rvd25yZXYueG1sRI/dagIx&#xA;FITvC75DOELvNKvQqqtRpD9Q8aJUfYDj5rhZ3ZwsSapbn94IQi+HmfmGmS1aW4sz+VA5VjDoZyCI&#xA;C6crLhXstp+9MYgQkTXWjknBHwVYzDtPM8y1u/APnTexFAnCIUcFJsYmlzIUhiyGvmuIk3dw3mJM&#xA;0pdSe7wkuK3lMMtepcWK04LBht4MFafNr1Ww8vv1aXAtjdzzyn/U3++TYI9KPXfb5RREpDb+hx/t&#xA;L61gOHoZwf1NegJyfgMAAP//AwBQSwECLQAUAAYACAAAACEA2+H2y+4AAACFAQAAEwAAAAAAAAAA&#xA;AAAAAAAAAAAAW0NvbnRlbnRfVHlwZXNdLnhtbFBLAQItABQABgAIAAAAIQBa9CxbvwAAABUBAAAL&#xA;AAAAAAAAAAAAAAAAAB8BAABfcmVscy8ucmVsc1BLAQItABQABgAIAAAAIQAuRlE/xQAAAN0AAAAP&#xA;AAAAAAAAAAAAAAAAAAcCAABkcnMvZG93bnJldi54bWxQSwUGAAAAAAMAAwC3AAAA+QIAAAAA&#xA;" strokeweight="1pt"/><v:rect id="Rectangle 969" o:spid="_x0000_s2010" style="position:absolute;left:8595;top:15651;width:2928;height:481;visibility:visible" o:gfxdata="UEsDBBQABgAIAAAAIQDb4fbL7gAAAIUBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbHyQz07DMAyH&#xA;70i8Q+QralM4IITa7kDhCAiNB7ASt43WOlEcyvb2pNu4IODoPz9/n1xv9vOkForiPDdwXVagiI23&#xA;jocG3rdPxR0oScgWJ8/UwIEENu3lRb09BBKV0ywNjCmFe63FjDSjlD4Q50nv44wpl3HQAc0OB9I3&#xA;VXWrjedEnIq03oC27qjHjympx31un0wiTQLq4bS4shrAECZnMGVTvbD9QSnOhDInjzsyuiBXWQP0&#xA;r4R18jfgnHvJr4nOknrFmJ5xzhraRtHWf3Kkpfz/yGo5S+H73hkquyhdjr3R8m2lj09svwAAAP//&#xA;AwBQSwMEFAAGAAgAAAAhAFr0LFu/AAAAFQEAAAsAAABfcmVscy8ucmVsc2zPwWrDMAwG4Ptg72B0&#xA;X5TuUMaI01uh19I+gLGVxCy2jGSy9e1nemrHjpL4P0nD4SetZiPRyNnCruvBUPYcYp4tXC/Htw8w&#xA;Wl0ObuVMFm6kcBhfX4Yzra62kC6xqGlKVgtLreUTUf1CyWnHhXKbTCzJ1VbKjMX5LzcTvvf9HuXR&#xA;gPHJNKdgQU5hB+ZyK23zHztFL6w81c5zQp6m6P9TMfB3PtPWFCczVQtB9N4U2rp2HOA44NMz4y8A&#xA;AAD//wMAUEsDBBQABgAIAAAAIQByUEDRwAAAAN0AAAAPAAAAZHJzL2Rvd25yZXYueG1sRE/LisIw&#xA;FN0P+A/hCu7GVPHR6RilDAhurQqzvDR32mpzU5OM1r83C8Hl4bxXm9604kbON5YVTMYJCOLS6oYr&#xA;BcfD9jMF4QOyxtYyKXiQh8168LHCTNs77+lWhErEEPYZKqhD6DIpfVmTQT+2HXHk/qwzGCJ0ldQO&#xA;7zHctHKaJAtpsOHYUGNHPzWVl+LfKMjzc3+6Fl+49TJN3ELPdJX/KjUa9vk3iEB9eItf7p1WMF3O&#xA;49z4Jj4BuX4CAAD//wMAUEsBAi0AFAAGAAgAAAAhANvh9svuAAAAhQEAABMAAAAAAAAAAAAAAAAA&#xA;AAAAAFtDb250ZW50X1R5cGVzXS54bWxQSwECLQAUAAYACAAAACEAWvQsW78AAAAVAQAACwAAAAAA&#xA;AAAAAAAAAAAfAQAAX3JlbHMvLnJlbHNQSwECLQAUAAYACAAAACEAclBA0cAAAADdAAAADwAAAAAA&#xA;AAAAAAAAAAAHAgAAZHJzL2Rvd25yZXYueG1sUEsFBgAAAAADAAMAtwAAAPQCAAAAAA==&#xA;" filled="f" stroked="f" strokeweight=".25pt"><v:textbox style="mso-next-textbox:#Rectangle 969" inset="1pt,1pt,1pt,1pt"><w:txbxContent><w:p w14:paraId="2BC7A300" w14:textId="5F7A75EF" w:rsidR="00B827FD" w:rsidRPr="00CB455E" w:rsidRDefault="00B827FD" w:rsidP="00B827FD"><w:pPr><w:ind w:right="57"/><w:jc w:val="center"/><w:rPr><w:lang w:val="ru-RU"/></w:rPr></w:pPr><w:r><w:rPr><w:i/><w:sz w:val="24"/><w:lang w:val="ru-RU"/></w:rPr><w:t>7</w:t></w:r><w:r><w:rPr><w:i/><w:sz w:val="24"/></w:rPr><w:t>4</w:t></w:r><w:r><w:rPr><w:i/><w:sz w:val="24"/><w:lang w:val="ru-RU"/></w:rPr><w:t>4</w:t></w:r><w:r><w:rPr><w:i/><w:sz w:val="24"/></w:rPr><w:t>1</w:t></w:r><w:r><w:rPr><w:i/><w:sz w:val="24"/><w:lang w:val="ru-RU"/></w:rPr><w:t>7006</w:t></w:r><w:r><w:rPr><w:i/><w:sz w:val="24"/></w:rPr><w:t>, 202</w:t></w:r><w:r><w:rPr><w:i/><w:sz w:val="24"/><w:lang w:val="ru-RU"/></w:rPr><w:t>1</w:t></w:r></w:p><w:p w14:paraId="24C4F15F" w14:textId="77777777" w:rsidR="00B827FD" w:rsidRPr="0042148F" w:rsidRDefault="00B827FD" w:rsidP="00B827FD"><w:pPr><w:pStyle w:val="a1"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="FFFFFF"/><w:szCs w:val="16"/><w:lang w:val="ru-RU"/></w:rPr></w:pPr></w:p></w:txbxContent></v:textbox></v:rect><v:rect id="Rectangle 970" o:spid="_x0000_s2011" style="position:absolute;left:8684;top:15204;width:540;height:330;visibility:visible" o:gfxdata="UEsDBBQABgAIAAAAIQDb4fbL7gAAAIUBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbHyQz07DMAyH&#xA;70i8Q+QralM4IITa7kDhCAiNB7ASt43WOlEcyvb2pNu4IODoPz9/n1xv9vOkForiPDdwXVagiI23&#xA;jocG3rdPxR0oScgWJ8/UwIEENu3lRb09BBKV0ywNjCmFe63FjDSjlD4Q50nv44wpl3HQAc0OB9I3&#xA;VXWrjedEnIq03oC27qjHjympx31un0wiTQLq4bS4shrAECZnMGVTvbD9QSnOhDInjzsyuiBXWQP0&#xA;r4R18jfgnHvJr4nOknrFmJ5xzhraRtHWf3Kkpfz/yGo5S+H73hkquyhdjr3R8m2lj09svwAAAP//&#xA;AwBQSwMEFAAGAAgAAAAhAFr0LFu/AAAAFQEAAAsAAABfcmVscy8ucmVsc2zPwWrDMAwG4Ptg72B0&#xA;X5TuUMaI01uh19I+gLGVxCy2jGSy9e1nemrHjpL4P0nD4SetZiPRyNnCruvBUPYcYp4tXC/Htw8w&#xA;Wl0ObuVMFm6kcBhfX4Yzra62kC6xqGlKVgtLreUTUf1CyWnHhXKbTCzJ1VbKjMX5LzcTvvf9HuXR&#xA;gPHJNKdgQU5hB+ZyK23zHztFL6w81c5zQp6m6P9TMfB3PtPWFCczVQtB9N4U2rp2HOA44NMz4y8A&#xA;AAD//wMAUEsDBBQABgAIAAAAIQAdHOVKxAAAAN0AAAAPAAAAZHJzL2Rvd25yZXYueG1sRI9Ba8JA&#xA;FITvQv/D8gq9mU2l2hhdJRSEXo0WenxkX5PY7Nt0d5uk/74rCB6HmfmG2e4n04mBnG8tK3hOUhDE&#xA;ldUt1wrOp8M8A+EDssbOMin4Iw/73cNsi7m2Ix9pKEMtIoR9jgqaEPpcSl81ZNAntieO3pd1BkOU&#xA;rpba4RjhppOLNF1Jgy3HhQZ7emuo+i5/jYKiuEwfP+UaD15mqVvpF10Xn0o9PU7FBkSgKdzDt/a7&#xA;VrB4Xa7h+iY+Abn7BwAA//8DAFBLAQItABQABgAIAAAAIQDb4fbL7gAAAIUBAAATAAAAAAAAAAAA&#xA;AAAAAAAAAABbQ29udGVudF9UeXBlc10ueG1sUEsBAi0AFAAGAAgAAAAhAFr0LFu/AAAAFQEAAAsA&#xA;AAAAAAAAAAAAAAAAHwEAAF9yZWxzLy5yZWxzUEsBAi0AFAAGAAgAAAAhAB0c5UrEAAAA3QAAAA8A&#xA;AAAAAAAAAAAAAAAABwIAAGRycy9kb3ducmV2LnhtbFBLBQYAAAAAAwADALcAAAD4AgAAAAA=&#xA;" filled="f" stroked="f" strokeweight=".25pt"><v:textbox style="mso-next-textbox:#Rectangle 970" inset="1pt,1pt,1pt,1pt"><w:txbxContent><w:p w14:paraId="326C9018" w14:textId="77777777" w:rsidR="00B827FD" w:rsidRPr="00B827FD" w:rsidRDefault="00B827FD" w:rsidP="00B827FD"><w:pPr><w:pStyle w:val="a1"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="ru-RU"/></w:rPr></w:pPr><w:r w:rsidRPr="00B827FD"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="ru-RU"/></w:rPr><w:t>У</w:t></w:r></w:p></w:txbxContent></v:textbox></v:rect><w10:wrap anchorx="margin" anchory="margin"/></v:group></w:pict></w:r></w:p><w:sectPr w:rsidR="00627F81" w:rsidRPr="00B827FD" w:rsidSect="00B827FD"><w:pgSz w:w="11906" w:h="16838" w:code="9"/><w:pgMar w:top="720" w:right="720" w:bottom="720" w:left="720" w:header="720" w:footer="720" w:gutter="0"/><w:cols w:space="720"/><w:docGrid w:linePitch="381"/></w:sectPr>'
[void]$d.Content.InsertXML($bodyXml)
Write-Output "Applied screenshot/title-block fixes."
